$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_A10 = @'
2025-12-30 00:49:42
'@
$ws.Range("A10").Value = $v_A10

$v_B10 = @'
gemini-2.5-flash
'@
$ws.Range("B10").Value = $v_B10

$v_C10 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document f...
'@
$ws.Range("C10").Value = $v_C10

$v_D10 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document fragment.
Input: The attached text from a contract.
CRITICAL INSTRUCTIONS:
1. **Assume Isolation with Common Sense**: Do NOT assume missing definitions exist in other documents. However, IGNORE common commercial lending terms typically defined in a base Credit Agreement (e.g., "Borrower", "Administrative Agent", "Lender", "Business Day", "Dollars", "GAAP", "Material Adverse Effect"). Only flag specific, deal-specific, or unusual capitalized terms that are undefined.
2. **Logic Check:** Check all math and logic tables.
3. **Drafting Errors:** Find any placeholders like "[__]" or blank lines that were forgotten.
4. **Exact Quotes:** For every error, you MUST provide the `exact_quote` from the text that contains the error. This is used for highlighting.
Output Format:
Return ONLY a valid JSON object with the following structure:
{
  "errors": [
    {
      "location": "Page 3, Section 2.1",
      "error": "Description of the error",
      "suggestion": "Suggested fix",
      "exact_quote": "Exact text substring to highlight in red"
    }
  ]
}
If no errors are found, return {"errors": []}.
--- CONTRACT TEXT BEGINS ---
--- [START OF PAGE 1] ---
 EX-4.1 3 d45101 1dex41.htm EX-4.1 
 Exhibit 4.1 
 CAPIT AL ONE PRIME AUT O RECEIV ABLES TRUST 2023-1 
 Class A-1 4.900% Auto Loan Asset Backed Notes 
 Class A-2 5.20% Auto Loan Asset Backed Notes 
 Class A-3 4.87% Auto Loan Asset Backed Notes 
 Class A-4 4.76% Auto Loan Asset Backed Notes 
 Class B 5.07% Auto Loan Asset Backed Notes 
 Class C 5.36% Auto Loan Asset Backed Notes 
 Class D 6.34% Auto Loan Asset Backed Notes 
 INDENTURE 
 Dated as of February 23, 2023 
 WILMINGT ON TRUST , NATIONAL ASSOCIA TION, 
 as the Indenture T rustee 
--- [START OF PAGE 2] ---
 CROSS REFERENCE T ABLE  1 
 TIA  Section  Indenture  Section 
 310  (a) (1)  6.11 
 (a) (2)  6.11 
 (a) (3)  6.10; 6.1 1 
 (a) (4)  N.A.  2 
 (a) (5)  6.11 
 (b)  6.8; 6.1 1 
 (c)  N.A. 
 311  (a)  6.12 
 (b)  6.12 
 (c)  N.A. 
 312  (a)  7.1 
 (b)  7.2 
 (c)  7.2 
 313  (a)  7.3 
 (b) (1)  7.3 
 (b) (2)  7.3 
 (c)  7.3 
 (d)  7.3 
 314  (a)  3.9 
 (b)  11.1, 3.6 
 (c) (1)  11.1 
 (c) (2)  11.1 
 (c) (3)  11.1 
 (d)  11.1 
 (e)  11.1 
 (f)  N.A. 
 315  (a)  6.1(b) 
 (b)  6.5 
 (c)  6.1(a) 
 (d)  6.1(c) 
 (e)  5.13 
 316  (a) (1) (A)  5.11 
 (a) (1) (B)  5.12 
 (a) (2)  N.A. 
 (b)  5.7 
 (c)  5.6(b) 
 317  (a) (1)  5.3(b) 
 (a) (2)  5.3(d) 
 (b)  3.3 
 318  (a)  11.7 
 1  Note: This Cross Reference T able shall not, for  any purpose, be deemed to be part of this Indenture. 
 2  N.A. means Not Applicable. 
 COPAR 2023-1 Indentur e 
--- [START OF PAGE 3] ---
 Page 
 ARTICLE I  DEFINITIONS AND INCORPORA TION BY REFERENCE  2 
 Section 1.1  Definitions  2 
 Section 1.2  Incorporation by Reference of T rust Indenture Act  2 
 Section 1.3  Other Interpretive Provisions  2 
 ARTICLE II  THE NOTES  3 
 Section 2.1  Form  3 
 Section 2.2  Execution, Authentication and Delivery  3 
 Section 2.3  Temporary Notes  4 
 Section 2.4  Registration of T ransfer and Exchange  4 
 Section 2.5  Mutilated, Destroyed, Lost or Stolen Notes  6 
 Section 2.6  Persons Deemed Owners  6 
 Section 2.7  Payment of Principal and Interest; Defaulted Interest  7 
 Section 2.8  Cancellation  7 
 Section 2.9  Release of Collateral  8 
 Section 2.10  Book-Entry Notes  8 
 Section 2.1 1  Notices to Clearing Agency  9 
 Section 2.12  Definitive Notes  9 
 Section 2.13  Authenticating Agents  10 
 Section 2.14  Paying Agent  10 
 Section 2.15  Tax and Accounting T reatment  11 
 Section 2.16  Certain T ransfer Restrictions on all Classes of the Notes  11 
 Section 2.17  Certain T ransfer Restrictions  12 
 Section 2.18  Transfer Restrictions on Certain Notes Upon a Sale of a Certificate  12 
 Section 2.19  Certain T ransfer Restrictions on the 144A Notes  13 
 ARTICLE III  COVENANTS  17 
 Section 3.1  Payment of Principal and Interest  17 
 Section 3.2  Maintenance of Of fice or Agency  17 
 Section 3.3  Money for Payments to Be Held in T rust  17 
 Section 3.4  Existence  19 
 Section 3.5  Protection of Collateral  19 
--- [START OF PAGE 4] ---
--- [START OF PAGE 5] ---
 Page 
 Section 3.6  Opinions as to Collateral  20 
 Section 3.7  Performance of Obligations; Servicing of Receivables  20 
 Section 3.8  Negative Covenants  21 
 Section 3.9  Annual Compliance Statement  22 
 Section 3.10  Restrictions on Certain Other Activities  23 
 Section 3.1 1  Restricted Payments  23 
 Section 3.12  Notice of Events of Default  23 
 Section 3.13  Further Instruments and Acts  23 
 Section 3.14  Compliance with Laws  23 
 Section 3.15  Removal of Administrator  23 
 Section 3.16  Perfection Representations, W arranties and Covenants  24 
 Section 3.17  Investment Company Act Representation  24 
 ARTICLE IV  SATISF ACTION AND DISCHARGE  24 
 Section 4.1  Satisfaction and Dischar ge of Indenture  24 
 Section 4.2  Application of T rust Money  25 
 Section 4.3  Repayment of Monies Held by Paying Agent  25 
 ARTICLE V  REMEDIES  25 
 Section 5.1  Events of Default  25 
 Section 5.2  Acceleration of Maturity; W aiver of Event of Default  26 
 Section 5.3  Collection of Indebtedness and Suits for Enforcement by the Indenture T rustee  27 
 Section 5.4  Remedies; Priorities  29 
 Section 5.5  Optional Preservation of the Collateral  32 
 Section 5.6  Limitation of Suits  32 
 Section 5.7  Rights of Noteholders to Receive Principal and Interest  33 
 Section 5.8  Restoration of Rights and Remedies  33 
 Section 5.9  Rights and Remedies Cumulative  33 
 Section 5.10  Delay or Omission Not a W aiver  33 
 Section 5.1 1  Control by Noteholders  34 
 Section 5.12  Waiver of Past Defaults  34 
 Section 5.13  Undertaking for Costs  35 
--- [START OF PAGE 6] ---
--- [START OF PAGE 7] ---
 Page 
 Section 5.14  Waiver of Stay or Extension Laws  35 
 Section 5.15  Action on Notes  35 
 Section 5.16  Performance and Enforcement of Certain Obligations  35 
 Section 5.17  Sale of Collateral  36 
 ARTICLE VI  THE INDENTURE TRUSTEE  37 
 Section 6.1  Duties of the Indenture T rustee  37 
 Section 6.2  Rights of the Indenture T rustee  38 
 Section 6.3  Individual Rights of the Indenture T rustee  40 
 Section 6.4  The Indenture T rustee’ s Disclaimer  40 
 Section 6.5  Notice of Defaults  40 
 Section 6.6  Reports by the Paying Agent  40 
 Section 6.7  Compensation and Indemnity  40 
 Section 6.8  Removal, Resignation and Replacement of the Indenture T rustee  41 
 Section 6.9  Successor Indenture T rustee by Mer ger  42 
 Section 6.10  Appointment of Co-Indenture T rustee or Separate Indenture T rustee  43 
 Section 6.1 1  Eligibility; Disqualification  44 
 Section 6.12  Preferential Collection of Claims Against the Issuer  44 
 Section 6.13  Representations and W arranties  44 
 ARTICLE VII  NOTEHOLDERS’ LISTS AND REPOR TS  45 
 Section 7.1  The Issuer to Furnish the Indenture T rustee Names and Addresses of Noteholders  45 
 Section 7.2  Preservation of Information; Communications to Noteholders  45 
 Section 7.3  Reports by the Indenture T rustee  45 
 Section 7.4  Statements to Certificateholders and Noteholders  45 
 Section 7.5  Noteholder Demand for Repurchase, Dispute Resolution  47 
 Section 7.6  Investor Action to Initiate an Asset Review  48 
 ARTICLE VIII  ACCOUNTS, DISBURSEMENTS AND RELEASES  49 
 Section 8.1  Collection of Money  49 
 Section 8.2  Trust Accounts  49 
--- [START OF PAGE 8] ---
 Page 
 Section 8.3  General Provisions Regarding Accounts  51 
 Section 8.4  Additional W ithdrawals and Deposits  54 
 Section 8.5  Distributions  54 
 Section 8.6  Release of Collateral  56 
 Section 8.7  Opinion of Counsel  56 
 ARTICLE IX  SUPPLEMENT AL INDENTURES  57 
 Section 9.1  Supplemental Indentures W ithout Consent of Noteholders  57 
 Section 9.2  Supplemental Indentures with Consent of Noteholders  58 
 Section 9.3  Execution of Supplemental Indentures  59 
 Section 9.4  Effect of Supplemental Indenture  60 
 Section 9.5  Conformity with T rust Indenture Act  60 
 Section 9.6  Reference in Notes to Supplemental Indentures  60 
 ARTICLE X  REDEMPTION OF NOTES  60 
 Section 10.1  Redemption  60 
 Section 10.2  Form of Redemption Notice  61 
 Section 10.3  Notes Payable on Redemption Date  61 
 ARTICLE XI  MISCELLANEOUS  62 
 Section 1 1.1  Compliance Certificates and Opinions, etc.  62 
 Section 1 1.2  Form of Documents Delivered to the Indenture T rustee  63 
 Section 1 1.3  Acts of Noteholders  64 
 Section 1 1.4  Notices  65 
 Section 1 1.5  Notices to Noteholders; W aiver  65 
 Section 1 1.6  Alternate Payment and Notice Provisions  65 
 Section 1 1.7  Conflict with T rust Indenture Act  65 
 Section 1 1.8  Information Requests  66 
 Section 1 1.9  Effect of Headings and T able of Contents  66 
 Section 1 1.10  Successors and Assigns  66 
 Section 1 1.11  Separability  66 
 Section 1 1.12  Benefits of Indenture  66 
 Section 1 1.13  Legal Holidays  66 
--- [START OF PAGE 9] ---
--- [START OF PAGE 10] ---
 Page 
 Section 1 1.14  GOVERNING LA W  66 
 Section 1 1.15  Counterparts  66 
 Section 1 1.16  Recording of Indenture  67 
 Section 1 1.17  Trust Obligation  67 
 Section 1 1.18  No Petition  67 
 Section 1 1.19  Submission to Jurisdiction; W aiver of Jury T rial  68 
 Section 1 1.20  Subordination of Claims  68 
 Section 1 1.21  U.S.A Patriot Act  69 
 Section 1 1.22  Beneficial Ownership  69 
 Section 1 1.23  Limitation of Liability  69 
 ARTICLE XII  COMPLIANCE WITH THE FDIC RULE  70 
 Section 12.1  Purpose  70 
 Section 12.2  Requirements of the FDIC Rule  71 
 Section 12.3  Performance  72 
 Section 12.4  Actions Upon Repudiation  73 
 Section 12.5  Notice  75 
 Section 12.6  Reservation of Rights  75 
 Schedule I  Perfection Representations, W arranties and Covenants 
 Exhibit A  Form of Notes 
--- [START OF PAGE 11] ---
 This INDENTURE, dated as of February 28, 2023 (as amended, supplemented, or otherwise modified and in ef fect from time to time, this 
 “  Indenture  ”), is between CAPIT AL ONE PRIME AUT O RECEIV ABLES  TRUST 2023-1, a Delaware statutory trust (the “  Issuer  ”),  and W ilmington 
 Trust, National Association, a national banking association, solely as trustee and not in its individual capacity (the “  Indenture T rustee  ”). 
 Each party agrees as follows for the benefit of the other party and the equal and ratable benefit of the Holders of the Issuer ’s Class A-1 
 4.900% Auto Loan Asset Backed Notes (the “  Class A-1  Notes  ”), Class A-2 5.20% Auto Loan Asset Backed Notes  (the “  Class A-2 Notes  ”), Class 
 A-3 4.87% Auto Loan Asset Backed Notes (the “  Class  A-3 Notes  ”), Class A-4 4.76% Auto Loan Asset Backed  Notes (the “  Class A-4 Notes  ”, and 
 together with the Class A-1 Notes, the Class A-2 Notes and the Class A-3 Notes, the “  Class A Notes  ”), Class  B 5.07% Auto Loan Asset Backed 
 Notes (the “  Class B Notes  ”), Class C 5.36% Auto Loan  Asset Backed Notes (the “  Class C Notes  ”) and Class  D 6.34% Auto Loan Asset Backed 
 Notes (the “  Class D Notes  ”, and together with the  Class A-1 Notes, the Class A-2 Notes, the Class A-3 Notes, the Class A-4 Notes, the Class B Notes 
 and the Class C Notes, the “  Notes  ”). 
 GRANTING CLAUSE 
 The Issuer , to secure the payment of principal of and interest on, and any other amounts owing in respect of, the Notes equally and ratably 
 without prejudice, priority or distinction except as set forth herein, and to secure compliance with the provisions of this Indenture, hereby Grants in trust 
 to the Indenture T rustee on the Closing Date, as trustee for the benefit of the Noteholders, all of the Issuer ’s right, title and interest, whether now owned 
 or hereafter acquired, in and to (i) the T rust Estate and (ii) all present and future claims, demands, causes and choses in action in respect of any or all of 
 the T rust Estate and all payments on or under and all proceeds of every kind and nature whatsoever in respect of any or all of the T rust Estate, including 
 all proceeds of the conversion, voluntary or involuntary , into cash or other liquid property , all cash proceeds, accounts, accounts receivable, notes, drafts, 
 acceptances, chattel paper , checks, deposit accounts, insurance proceeds, condemnation awards, rights to payment of any and every kind and other forms 
 of obligations and receivables, instruments, securities, financial assets and other property which at any time constitute all or part of or are included in the 
 proceeds of any of the T rust Estate (collectively , the “  Collateral  ”). 
 The Indenture T rustee, on behalf of the Noteholders, acknowledges the foregoing Grant, accepts the trusts under this Indenture and 
 agrees to perform its duties required in this Indenture in accordance with the provisions of this Indenture. 
 The foregoing Grant is made in trust to secure (i) the payment of principal of and interest on, and any other amounts owing in respect of, 
 the Notes, equally and ratably without prejudice, priority or distinction except as set forth herein and (ii) compliance with the provisions of this 
 Indenture, all as provided in this Indenture. 
 Without limiting the foregoing Grant, any Receivable repurchased or purchased by (a) the Servicer pursuant to  Section 3.6  of the Servicing 
 Agreement or (b) by the Bank pursuant to  Section 3.4  of the Purchase Agreement shall be deemed to be automatically released from the lien of this 
 Indenture without any action being taken by the Indenture T rustee upon payment by the applicable purchaser of the related Repurchase Price for such 
 Repurchased Receivable. 
 COPAR 2023-1 Indentur e 
--- [START OF PAGE 12] ---
 ARTICLE I 
 DEFINITIONS AND INCORPORA TION BY REFERENCE 
 SECTION 1.1  Definitions  . Capitalized terms are used  in this Indenture as defined in  A  p  p  endix A  to the  Sale Agreement, dated as of the 
 date hereof (as amended, supplemented, or otherwise modified and in ef fect from time to time, the “  Sale  Agreement  ”), between the Issuer and Capital 
 One Auto Receivables, LLC, which also contains rules as to usage that are applicable herein. 
 SECTION 1.2  Incorporation by Reference of T rust Indenture  Act  . Whenever this Indenture refers to a provision  of the TIA, the provision 
 is incorporated by reference in and made a part of this Indenture. The following TIA terms used in this Indenture have the following meanings: 
 “  Commission  ” means the Securities and Exchange Commission. 
 “  indenture securities  ” means the Notes. 
 “  indenture securi  ty  holder  ” means a Noteholder . 
 “  indenture to be qualified  ” means this Indenture. 
 “  indenture trustee  ” or “  institutional trustee  ” means  the Indenture T rustee. 
 “  obl  ig  or  ” on the indenture securities means the Issuer  and any other obligor on the indenture securities. 
 All other TIA terms used in this Indenture that are defined by the TIA, defined by TIA reference to another statute or defined by 
 Commission rule have the meaning assigned to them by such definitions. 
 SECTION 1.3  Other Interpretive Provisions  . All terms  defined in this Indenture shall have the defined meanings when used in any 
 certificate or other document delivered pursuant hereto unless otherwise defined therein. For purposes of this Indenture and all such certificates and 
 other documents, unless the context otherwise requires: (a) accounting terms not otherwise defined in this Indenture, and accounting terms partly 
 defined in this Indenture to the extent not defined, shall have the respective meanings given to them under GAAP (  provided  , that, to the extent that the 
 definitions in this Indenture and GAAP conflict, the definitions in this Indenture shall control); (b) terms defined in Article 9 of the UCC as in ef fect in 
 the relevant jurisdiction and not otherwise defined in this Indenture are used as defined in that Article; (c) the words “hereof,” “herein” and “hereunder” 
 and words of similar import refer to this Indenture as a whole and not to any particular provision of this Indenture; (d) references to any Article, 
 Section, Schedule, Appendix or Exhibit are references to Articles, Sections, Schedules, Appendices and Exhibits in or to this Indenture and references 
 to any paragraph, subsection, clause or other subdivision within any Section or definition refer to such paragraph, subsection, clause or other 
--- CONTRACT TEXT ENDS ---
'@
$ws.Range("D10").Value = $v_D10

$v_E10 = @'
{
  "errors": [
    {
      "location": "Page 11",
      "error": "The term \"Holders\" is capitalized but not explicitly defined. While \"Noteholders\" is used frequently and its meaning is implied through cross-references to the TIA, inconsistent capitalization without a clear definition can lead to ambiguity.",
      "suggestion": "Either define \"Holders\" or consistently use \"Noteholders\" (assuming \"Noteholders\" is a defined term in the Sale Agreement as per Section 1.1).",
      "exact_quote": "the Holders of the Issuer’s Class A-1"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The verb \"Grants\" is capitalized as if it were a defined term, but it is used grammatically as an action verb.",
      "suggestion": "Decapitalize \"Grants\" to \"grants\" unless it is specifically intended to be a defined legal term, in which case it should be explicitly defined and used consistently (e.g., \"makes a Grant\").",
      "exact_quote": "hereby Grants in trust"
    },
    {
      "location": "Page 12, Section 1.3",
      "error": "The last sentence of Section 1.3 ends abruptly and is incomplete.",
      "suggestion": "Complete the sentence or remove the incomplete phrase \"clause or other\" to ensure proper grammar and meaning.",
      "exact_quote": "clause or other"
    },
    {
      "location": "Page 12",
      "error": "The footer on this page (internal page 2 of the document, as per the Table of Contents) is missing its page number.",
      "suggestion": "Add the internal page number '2' to the footer of this page, consistent with the Table of Contents and standard document formatting.",
      "exact_quote": "COPAR 2023-1 Indentur e"
    }
  ]
}
'@
$ws.Range("E10").Value = $v_E10

$v_A11 = @'
2025-12-30 00:55:13
'@
$ws.Range("A11").Value = $v_A11

$v_B11 = @'
gemini-2.5-flash
'@
$ws.Range("B11").Value = $v_B11

$v_C11 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document f...
'@
$ws.Range("C11").Value = $v_C11

$v_D11 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document fragment.
Input: The attached text from a contract.
CRITICAL INSTRUCTIONS:
1. **Assume Isolation with Common Sense**: Do NOT assume missing definitions exist in other documents. However, IGNORE common commercial lending terms typically defined in a base Credit Agreement (e.g., "Borrower", "Administrative Agent", "Lender", "Business Day", "Dollars", "GAAP", "Material Adverse Effect"). Only flag specific, deal-specific, or unusual capitalized terms that are undefined.
2. **Logic Check:** Check all math and logic tables.
3. **Drafting Errors:** Find any placeholders like "[__]" or blank lines that were forgotten.
4. **Exact Quotes:** For every error, you MUST provide the `exact_quote` from the text that contains the error. This is used for highlighting.
Output Format:
Return ONLY a valid JSON object with the following structure:
{
  "errors": [
    {
      "location": "Page 3, Section 2.1",
      "error": "Description of the error",
      "suggestion": "Suggested fix",
      "exact_quote": "Exact text substring to highlight in red"
    }
  ]
}
If no errors are found, return {"errors": []}.
--- CONTRACT TEXT BEGINS ---
--- [START OF PAGE 1] ---
EX-10 2 tmb-20250626xex10.htm EX-10
FIRST AMENDMENT TO AMENDED AND RESTATED CREDIT
AGREEMENT
FIRST AMENDMENT TO AMENDED AND RESTATED CREDIT AGREEMENT
(this “Agreement”), dated as of June 26, 2025 (and effective as set forth below), among
MAYVILLE ENGINEERING COMPANY, INC., a Wisconsin corporation (the “Borrower”),
the Subsidiary Guarantors (as defined in the Existing Credit Agreement referred to below)
party hereto, the Incremental Lenders (as defined below), the Lenders (as defined below)
party hereto (the “Consenting Lenders”), and WELLS FARGO BANK, NATIONAL
ASSOCIATION, as administrative agent (the “Administrative Agent”). Unless otherwise
indicated, all capitalized terms used herein and not otherwise defined herein shall have the
respective meanings provided such terms in the Credit Agreement referred to below.
W I T N E S S E T H:
WHEREAS, the Borrower, the lenders party thereto (the “Lenders”) and the
Administrative Agent have entered into that certain Amended and Restated Credit Agreement,
dated as of June 28, 2023 (as amended, restated, supplemented or otherwise modified prior to
the date hereof, the “Existing Credit Agreement” and as hereby amended and as from time to
time further amended, modified, supplemented, restated, or amended and restated, the “Credit
Agreement”);
WHEREAS, the Borrower has requested an Incremental Revolving Credit Increase in
the aggregate principal amount of $100,000,000, in accordance with Section 5.13 of the
Existing Credit Agreement (the “First Amendment Revolving Credit Increase”);
WHEREAS, subject to the terms of this Agreement, certain of the Revolving Credit
Lenders identified on the signature pages hereto (each an “Incremental Lender”) are severally
willing to provide a portion of the First Amendment Revolving Credit Increase; and
WHEREAS, the Borrower has also requested, and subject to the terms and conditions
set forth herein, the Administrative Agent and the Consenting Lenders have agreed to further
amend the Existing Credit Agreement as more specifically set forth herein.
NOW, THEREFORE, in consideration of the foregoing and for other good and
valuable consideration, the receipt and sufficiency of which are hereby acknowledged, it is
agreed as follows:
SECTION 1.Amendments to Existing Credit Agreement. Effective as of the First
Amendment Effective Date (as defined below) and subject to the terms and conditions set
forth herein and in reliance upon the representations and warranties set forth herein:
(a)Schedule 1.1 (Commitments and Commitment Percentages) to the Existing
Credit Agreement is hereby amended and restated in the form attached hereto
as Annex A.
(b)the body of Existing Credit Agreement is hereby amended as follows: (i) to
delete red or green stricken text (indicated textually in the same manner as the
following examples: stricken text and stricken text) and (ii) to add the blue or
green double-underlined text (indicated textually in the same manner as the
following examples: double-underlined text and double-underlined text), in
each case, as set forth in the conformed copy attached hereto as Annex B.
--- [START OF PAGE 2] ---
1
--- [START OF PAGE 3] ---
SECTION 2.Revolving Credit Facility Increase; Reallocation.
(a)Each Incremental Lender severally agrees that its Revolving Credit
Commitment after giving effect to the First Amendment Revolving Credit
Increase shall be as set forth opposite such Lender’s name on Annex A hereto.
(b)The parties hereto hereby agree that (i) this Agreement constitutes a request
for an Incremental Revolving Credit Increase pursuant to Section 5.13 of the
Credit Agreement, (ii) the parties hereto hereby waive any prior notice
required thereby, and (iii) this Agreement shall be deemed to be an “Lender
Joinder Agreement” in accordance with Section 5.13(h) of the Credit
Agreement.
(c)Pursuant to Section 5.13 of the Credit Agreement, on the First Amendment
Effective Date (as defined below) and after giving effect to this Agreement,
the outstanding Revolving Credit Loans and Revolving Credit Commitment
Percentages of Swingline Loans and L/C Obligations will be reallocated by
the Administrative Agent among the Revolving Credit Lenders (including the
Incremental Lenders providing the First Amendment Revolving Credit
Increase) in accordance with their revised Revolving Credit Commitment
Percentages (and the Revolving Credit Lenders (including the Incremental
Lenders providing the First Amendment Revolving Credit Increase) agree to
make all payments and adjustments necessary to effect such reallocation and
the Borrower shall pay any and all costs required pursuant to Section 5.9 of
the Credit Agreement in connection with such reallocation as if such
reallocation were a repayment).
SECTION 3.Conditions of Effectiveness of this Agreement.  This Agreement shall
become effective on the date when the following conditions shall have been satisfied or
waived (such date, the “First Amendment Effective Date”):
(a) The Administrative Agent’s receipt of the following, each in form and
substance reasonably satisfactory to the Administrative Agent:
(i) this Agreement, duly executed by each of the Credit Parties existing
as of the First Amendment Effective Date, the Administrative Agent, each
Incremental Lender and each of the Consenting Lenders;
(ii) a certificate of a Responsible Officer (or the secretary or assistant
secretary) of each Credit Party certifying as to the incumbency and genuineness of the
signature of each officer of such Credit Party executing Loan Documents to which it
is a party and certifying that attached thereto is a true, correct and complete copy of
(A) the articles or certificate of incorporation or formation (or equivalent), as
applicable, of such Credit Party and all amendments thereto, certified as of a recent
date by the appropriate Governmental Authority in its jurisdiction of incorporation,
organization or formation (or equivalent), as applicable, (B) the bylaws or other
governing document of such Credit Party as in effect on the First Amendment
Effective Date, (C) resolutions duly adopted by the board of directors (or other
governing body) of such Credit Party authorizing and approving the transactions
contemplated hereunder and the execution, delivery and performance of this
Amendment, the Credit Agreement and the other Loan Documents to which it is a
party, and (D) certificates as of a recent date of the good standing of each Credit Party
--- [START OF PAGE 4] ---
2under the laws of its jurisdiction of incorporation, organization or formation (or
equivalent), as applicable;
(iii) an Officer’s Compliance Certificate demonstrating, in form and
substance reasonably satisfactory to the Administrative Agent, that the Borrower is in
compliance with the
--- [START OF PAGE 5] ---
financial covenants set forth in Section 9.15 of the Credit Agreement as of the First
Amendment Effective Date based on the financial statements most recently delivered
pursuant to Section 8.1(a) or 8.1(b) of the Existing Credit Agreement, as applicable,
both before and after giving effect (on a pro forma basis) to (x) any Incremental
Revolving Credit Commitment, and (y) the making of the First Amendment
Revolving Credit Increase pursuant thereto (with such Incremental Revolving Credit
Commitment being deemed to be fully funded);
(iv) a certificate from the chief financial officer of the Borrower with
respect to solvency (on a consolidated basis) of the Borrower and its Subsidiaries as
of the First Amendment Effective Date after giving effect to the First Amendment
Revolving Credit Increase and the other amendments contemplated by this
Agreement;
(v) an opinion from counsel to the Credit Parties, in form and substance
reasonably satisfactory to the Administrative Agent.
(b) The representations and warranties set forth in Section 3 shall be true and
correct.
(c) No Default or Event of Default shall exist on the First Amendment Effective
Date immediately prior to or after giving effect to this Agreement.
(d) The Administrative Agent shall have received the results of a Lien search, in
form and substance reasonably satisfactory thereto, made against the Credit Parties under the
Uniform Commercial Code (or applicable judicial docket) as in effect in each jurisdiction in
which filings or recordations under the Uniform Commercial Code should be made to
evidence or perfect security interests in all assets of such Credit Party, indicating among other
things that the assets of each such Credit Party are free and clear of any Lien (except for
Permitted Liens).
(e) All governmental, shareholder and third party consents and approvals
necessary (or any other material consents as determined in the reasonable discretion of the
Administrative Agent) in connection with the transactions contemplated by this Agreement
shall have been obtained and shall be in full force and effect.
(f) The Borrower and each of the Subsidiary Guarantors shall have provided to
the Administrative Agent and the Lenders, at least five (5) Business Days prior to the First
Amendment Effective Date, a Beneficial Ownership Certification under 31 C.F.R. 1010.230
and all documentation and other information required by regulatory authorities under
applicable “know your customer”, anti-money-laundering, beneficial ownership and other
similar rules and regulations, including, without limitation, the PATRIOT Act.
(g) The Borrower shall have paid all fees and expenses required to be paid on the
First Amendment Effective Date.
For purposes of determining compliance with the conditions specified in this Section 3, the
Administrative Agent, each Incremental Lender and each Consenting Lender that has signed
this Agreement shall be deemed to have consented to, approved or accepted or to be satisfied
with, each document or other matter required thereunder to be consented to or approved by or
acceptable or satisfactory to an Incremental Lender and/or a Consenting Lender, as applicable,
unless the Administrative Agent shall have received notice from such Incremental Lender
--- [START OF PAGE 6] ---
3and/or Consenting Lender, as applicable, prior to the proposed First Amendment Effective
Date specifying its objection thereto.
--- [START OF PAGE 7] ---
SECTION 4.Representations and Warranties.  To induce the Administrative Agent,
the Incremental Lenders party hereto and the Consenting Lenders party hereto to enter into
this Agreement, each Credit Party represents and warrants to the Administrative Agent, the
Incremental Lenders and the Consenting Lenders on and as of the First Amendment Effective
Date that, in each case:
(a) each of the representations and warranties contained in Article VII of the
Credit Agreement are true and correct in all material respects, except to the extent any such
representation and warranty is qualified by materiality or reference to Material Adverse
Effect, in which case, such representation and warranty is true, correct and complete in all
respects, with the same effect as if made on and as of such date (except for any such
representation and warranty that by its terms is made only as of an earlier date, which
representation and warranty remains true and correct as of such earlier date);
(b) no Default or Event of Default has occurred and is continuing immediately
prior to or after giving effect to this Agreement;
(c) it has the right, power and authority and has taken all necessary corporate and
other action to authorize the execution, delivery and performance of this Agreement in
accordance with their respective terms; and
(d) this Agreement has been duly executed and delivered by the duly authorized
officers, managers, managing members or other authorized person of each Credit Party that is
a party thereto, and each such document constitutes the legal, valid and binding obligation of
each Credit Party, enforceable in accordance with its terms, except as such enforceability may
be limited by bankruptcy, insolvency, reorganization, moratorium or similar state or federal
Debtor Relief Laws from time to time in effect which affect the enforcement of creditors’
rights in general and the availability of equitable remedies.
SECTION 5.Credit Party Acknowledgements, Agreements and Confirmations.  By
its execution of this Agreement, each Credit Party hereby (a) consents to this Agreement and
agrees that the transactions contemplated by this Agreement shall not limit or diminish the
obligations of such Credit Party, or release such Credit Party from any obligations, under any
of the Loan Documents to which it is a party, (b) confirms and reaffirms its obligations under
each of the Loan Documents to which it is a party and (c) agrees that each of the Loan
Documents to which it is a party remain in full force and effect and are hereby ratified and
confirmed.
SECTION 6.Reference to and Effect on the Credit Agreement and the Loan
Documents.
(a) On and after the First Amendment Effective Date, each reference in the Credit
Agreement to “this Agreement,” “herein,” “hereto”, “hereof” and “hereunder” or words of like
import referring to the Credit Agreement, and each reference in the Notes and each of the
other Loan Documents to “the Credit Agreement”, “thereunder”, “thereof” or words of like
import referring to the Credit Agreement, shall mean and be a reference to the Credit
Agreement, as amended and supplemented by this Agreement.
(b) The Credit Agreement and each of the other Loan Documents, as specifically
amended and supplemented by this Agreement, are and shall continue to be in full force and
effect and are hereby in all respects ratified and confirmed.
(c) The execution, delivery and effectiveness of this Agreement shall not, except
as expressly provided herein, operate as a waiver of any right, power or remedy of any Lender
or the Administrative Agent under any of the Loan Documents, nor constitute a waiver of any
provision of any of the Loan Documents.  Without limiting the generality of the foregoing, the
--- [START OF PAGE 8] ---
4Security Documents in effect immediately prior to the date hereof and all of the Collateral
described therein in existence immediately prior to the date hereof do and shall continue to
secure the payment of all Secured Obligations.  
--- [START OF PAGE 9] ---
5(d) This Agreement is the entire agreement, and supersedes any prior agreements
and contemporaneous oral agreements, of the parties concerning its subject matter. This
Agreement shall constitute a “Loan Document” for all purposes of the Credit Agreement and
the other Loan Documents.
(e) This Agreement shall be binding on and inure to the benefit of the parties
hereto and their successors and permitted assigns.
SECTION 7.Costs and Expenses.  The Borrower hereby reconfirms its obligations
pursuant to Section 12.3 of the Credit Agreement to pay and reimburse the Administrative
Agent for all reasonable fees, costs and expenses in accordance with the terms thereof.
SECTION 8.Governing Law; Jury Trial Waiver.  This Agreement and any claim,
controversy, dispute or cause of action (whether in contract or tort or otherwise) based upon,
arising out of or relating to this Agreement and the transactions contemplated hereby shall be
governed by, and construed in accordance with, the law of the State of New York.  EACH
PARTY HERETO HEREBY IRREVOCABLY WAIVES, TO THE FULLEST EXTENT
PERMITTED BY APPLICABLE LAW, ANY RIGHT IT MAY HAVE TO A TRIAL BY JURY
IN ANY LEGAL PROCEEDING DIRECTLY OR INDIRECTLY ARISING OUT OF OR
RELATING TO THIS AGREEMENT OR THE TRANSACTIONS CONTEMPLATED
HEREBY (WHETHER BASED ON CONTRACT, TORT OR ANY OTHER THEORY).  EACH
PARTY HERETO (A) CERTIFIES THAT NO REPRESENTATIVE, AGENT OR ATTORNEY
OF ANY OTHER PERSON HAS REPRESENTED, EXPRESSLY OR OTHERWISE, THAT
SUCH OTHER PERSON WOULD NOT, IN THE EVENT OF LITIGATION, SEEK TO
ENFORCE THE FOREGOING WAIVER AND (B) ACKNOWLEDGES THAT IT AND THE
OTHER PARTIES HERETO HAVE BEEN INDUCED TO ENTER INTO THIS AGREEMENT
BY, AMONG OTHER THINGS, THE MUTUAL WAIVERS AND CERTIFICATIONS IN THIS
SECTION.
SECTION 9.Counterparts; Electronic Execution. This Agreement may be executed
in counterparts (and by different parties hereto in different counterparts), each of which shall
constitute an original, but all of which when taken together shall constitute a single contract.
Delivery of an executed counterpart of a signature page of this Agreement by facsimile or in
electronic (i.e., “pdf” or “tif”) format shall be effective as delivery of a manually executed
counterpart of this Agreement. This Agreement may be executed by Electronic Signatures or
in the form of an Electronic Record pursuant to, and in accordance with, the provisions of
Section 12.16 of the Credit Agreement.
[Signature pages to follow]
--- [START OF PAGE 10] ---
Mayville Engineering Company, Inc.IN WITNESS WHEREOF, the parties hereto have caused their duly authorized
officers to execute and deliver this Agreement as of the date first above written.
BORROWER:
MAYVILLE ENGINEERING COMPANY,
INC.
SUBSIDIARY GUARANTORS:
CENTER MANUFACTURING
HOLDINGS, INC.
CENTER MANUFACTURING, INC.
CENTER – MOELLER PRODUCTS LLC
DEFIANCE METAL PRODUCTS CO.
DEFIANCE METAL PRODUCTS OF
ARKANSAS, INC. 
 By:
Name:
Title:
By:
Name:
Title:
By:
Name:
Title:
By:
Name:
Title:
By:
Name:
Title:
By:
Name:
Title:
--- [START OF PAGE 11] ---
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature Page
--- [START OF PAGE 12] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageDEFIANCE METAL PRODUCTS OF PA.,
INC.
MID-STATES ALUMINUM, LLC 
 By:
Name:
Title:
By:
Name:
Title:
--- [START OF PAGE 13] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageWELLS FARGO BANK, NATIONAL
ASSOCIATION, as Administrative Agent,
Issuing Lender, Swingline Lender,
Consenting Lender and an Incremental
Lender
By:
___________________________________
Name:
Title:  
--- [START OF PAGE 14] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageJPMORGAN CHASE BANK, N.A., as a
Consenting Lender and an Incremental
Lender
By:
___________________________________
Name:
Title:  
--- [START OF PAGE 15] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageTD BANK, N.A., as a Consenting Lender
and an Incremental Lender
By:
___________________________________
Name:
Title:  
--- [START OF PAGE 16] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageCITIBANK, N.A., as a Consenting Lender
By:
___________________________________
Name:
Title:  
--- [START OF PAGE 17] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageHANMI BANK, as a Consenting Lender and
an Incremental Lender
By:
___________________________________
Name:
Title:  
--- [START OF PAGE 18] ---
Mayville Engineering Company, Inc.
First Amendment to A&R Credit Agreement (Revolver Increase)
Signature PageASSOCIATED BANK, N.A., as a
Consenting Lender and an Incremental
Lender
By:
___________________________________
Name:
Title:  
--- [START OF PAGE 19] ---
Annex A
REVOLVING CREDIT COMMITMENTS AND REVOLVING CREDIT COMMITMENT
PERCENTAGES
Lender Revolving
Credit
CommitmentRevolving
Credit
Commitment
Percentage
Wells Fargo Bank, National
Association 
$108,000,000.00 30.857142858%
JPMorgan Chase Bank, N.A.  $90,000,000.00 25.714285714%
TD Bank, N.A.  $77,000,000.00 22.000000000%
Associated Bank, N.A.  $30,000,000.00 8.571428571%
Hanmi Bank  $25,000,000.00 7.142857143%
Citibank, N.A.  $20,000,000.00 5.714285714%
Total $350,000,000.00100.000000000%
--- [START OF PAGE 20] ---
ANNEX B
Published CUSIP Number: 57859VAP8
$350,000,000.00
AMENDED AND RESTATED CREDIT AGREEMENT
dated as of June 28, 2023
(as amended by that certain First Amendment to Amended and Restated Credit Agreement
dated as of June 26, 2025)
by and among
MAYVILLE ENGINEERING COMPANY, INC.
as Borrower,
the Lenders referred to herein,
as Lenders,
and
WELLS FARGO BANK, NATIONAL ASSOCIATION,
as Administrative Agent,
Swingline Lender and Issuing Lender
and
WELLS FARGO SECURITIES, LLC,
as Joint Lead Arranger and Joint Bookrunner
and
J.P. MORGAN SECURITIES LLC,
as Joint Lead Arranger and Joint Bookrunner
and
TD BANK, N.A.,
as Joint Lead Arranger and Joint Bookrunner 
--- [START OF PAGE 21] ---
--- [START OF PAGE 22] ---
i 
 ARTICLE IDEFINITIONS 1
SECTION 1.1Definitions 1
SECTION 1.2Other Definitions and Provisions 31
SECTION 1.3Accounting Terms 31
SECTION 1.4UCC Terms 32
SECTION 1.5Rounding 32
SECTION 1.6References to Agreement and Laws 32
SECTION 1.7Times of Day 32
SECTION 1.8Letter of Credit Amounts 32
SECTION 1.9Guarantees/Earn-Outs 33
SECTION 1.10Covenant Compliance Generally 33
SECTION 1.11Rates 33
SECTION 1.12Divisions 33
ARTICLE IIREVOLVING CREDIT FACILITY 34
SECTION 2.1Revolving Credit Loans 34
SECTION 2.2Swingline Loans 34
SECTION 2.3Procedure for Advances of Revolving Credit Loans and Swingline
Loans 36
SECTION 2.4Repayment and Prepayment of Revolving Credit Loans and Swingline
Loans 37
SECTION 2.5Permanent Reduction of the Revolving Credit Commitment 37
SECTION 2.6Termination of Revolving Credit Facility 38
ARTICLE IIILETTER OF CREDIT FACILITY 38
SECTION 3.1L/C Facility 38
SECTION 3.2Procedure for Issuance of Letters of Credit 39
SECTION 3.3Commissions and Other Charges 39
SECTION 3.4L/C Participations 40
SECTION 3.5Reimbursement Obligation of the Borrower 41
SECTION 3.6Obligations Absolute 41
SECTION 3.7Effect of Letter of Credit Application 41
SECTION 3.8Resignation of Issuing Lenders 42
SECTION 3.9Reporting of Letter of Credit Information and L/C Commitment 42
SECTION 3.10Letters of Credit Issued for Subsidiaries 42
--- [START OF PAGE 23] ---
ii 
 ARTICLE IVReserved 43
ARTICLE VGENERAL LOAN PROVISIONS 43
SECTION 5.1Interest 43
SECTION 5.2Notice and Manner of Conversion or Continuation of Loans 44
SECTION 5.3Fees 44
SECTION 5.4Manner of Payment 45
SECTION 5.5Evidence of Indebtedness 45
SECTION 5.6Sharing of Payments by Lenders 46
SECTION 5.7Administrative Agent’s Clawback 46
SECTION 5.8Changed Circumstances 47
SECTION 5.9Indemnity 50
SECTION 5.10Increased Costs 50
SECTION 5.11Taxes 51
SECTION 5.12Mitigation Obligations; Replacement of Lenders 55
SECTION 5.13Incremental Loans 56
SECTION 5.14Cash Collateral 58
SECTION 5.15Defaulting Lenders 58
ARTICLE VICONDITIONS OF CLOSING AND BORROWING 60
SECTION 6.1Conditions to Closing and Initial Extensions of Credit 60
SECTION 6.2Conditions to All Extensions of Credit 64
ARTICLE VIIREPRESENTATIONS AND WARRANTIES OF THE CREDIT PARTIES64
SECTION 7.1Organization; Power; Qualification 64
SECTION 7.2Ownership 65
SECTION 7.3Authorization; Enforceability 65
SECTION 7.4Compliance of Agreement, Loan Documents and Borrowing with
Laws, Etc 65
SECTION 7.5Compliance with Law; Governmental Approvals 65
SECTION 7.6Tax Returns and Payments 66
SECTION 7.7Intellectual Property Matters 66
SECTION 7.8Environmental Matters 66
SECTION 7.9Employee Benefit Matters 67
SECTION 7.10Margin Stock 68
SECTION 7.11Government Regulation 68
SECTION 7.12Material Contracts 68
SECTION 7.13Employee Relations 68
--- [START OF PAGE 24] ---
--- [START OF PAGE 25] ---
iii 
 SECTION 7.14Burdensome Provisions 69
SECTION 7.15Financial Statements 69
SECTION 7.16No Material Adverse Change 69
SECTION 7.17Solvency 69
SECTION 7.18Title to Properties 69
SECTION 7.19Litigation 69
SECTION 7.20Anti-Corruption Laws; Anti-Money Laundering Laws; Sanctions and
Outbound Investment Rules 69
SECTION 7.21Absence of Defaults 70
SECTION 7.22ESOP 70
SECTION 7.23Senior Indebtedness Status 71
SECTION 7.24Disclosure 71
ARTICLE VIIIAFFIRMATIVE COVENANTS 72
SECTION 8.1Financial Statements and Budgets 72
SECTION 8.2Certificates; Other Reports 73
SECTION 8.3Notice of Litigation and Other Matters 73
SECTION 8.4Preservation of Corporate Existence and Related Matters 75
SECTION 8.5Maintenance of Property and Licenses 75
SECTION 8.6Insurance 75
SECTION 8.7Accounting Methods and Financial Records 76
SECTION 8.8Payment of Taxes and Other Obligations 76
SECTION 8.9Compliance with Laws and Approvals 76
SECTION 8.10Environmental Laws 76
SECTION 8.11Compliance with ERISA 76
SECTION 8.12Compliance with Material Contracts 76
SECTION 8.13Visits and Inspections 76
SECTION 8.14Additional Subsidiaries 77
SECTION 8.15Depository Bank 77
SECTION 8.16Use of Proceeds 78
SECTION 8.17ESOP 78
SECTION 8.18Compliance with Anti-Corruption Laws; Beneficial Ownership
Regulation, Anti-Money Laundering Laws and Sanctions 78
SECTION 8.19Further Assurances 78
SECTION 8.20Post-Closing Matters 79
SECTION 8.21Lender Meetings 79
--- [START OF PAGE 26] ---
--- [START OF PAGE 27] ---
iv 
 ARTICLE IXNEGATIVE COVENANTS 79
SECTION 9.1Indebtedness 79
SECTION 9.2Liens 80
SECTION 9.3Investments 81
SECTION 9.4Fundamental Changes 82
SECTION 9.5Asset Dispositions 83
SECTION 9.6Restricted Payments 84
SECTION 9.7Transactions with Affiliates 84
SECTION 9.8Accounting Changes; Organizational Documents 84
SECTION 9.9Payments and Modifications of Subordinated Indebtedness 85
SECTION 9.10No Further Negative Pledges; Restrictive Agreements 85
SECTION 9.11ESOP Matters 86
SECTION 9.12Sale of Accounts 86
SECTION 9.13Sale Leasebacks 86
SECTION 9.14Capital Expenditures 86
SECTION 9.15Financial Covenants 86
SECTION 9.16No Amendment to ESOP Documentation 87
SECTION 9.17Limitations Regarding Outbound Investment Rules 87
ARTICLE XDEFAULT AND REMEDIES 87
SECTION 10.1Events of Default 87
SECTION 10.2Remedies 89
SECTION 10.3Rights and Remedies Cumulative; Non-Waiver; etc 90
SECTION 10.4Crediting of Payments and Proceeds 91
SECTION 10.5Administrative Agent May File Proofs of Claim 91
SECTION 10.6Credit Bidding 92
ARTICLE XITHE ADMINISTRATIVE AGENT 92
SECTION 11.1Appointment and Authority 92
SECTION 11.2Rights as a Lender 93
SECTION 11.3Exculpatory Provisions 93
SECTION 11.4Reliance by the Administrative Agent 94
SECTION 11.5Delegation of Duties 94
SECTION 11.6Resignation of Administrative Agent 95
SECTION 11.7Non-Reliance on Administrative Agent and Other Lenders 96
SECTION 11.8No Other Duties, Etc 96
SECTION 11.9Collateral and Guaranty Matters 96
--- [START OF PAGE 28] ---
--- [START OF PAGE 29] ---
v 
 SECTION 11.10Secured Hedge Agreements and Secured Cash Management
Agreements 97
SECTION 11.11[Reserved.] 97
SECTION 11.12Erroneous Payments 97
ARTICLE XIIMISCELLANEOUS 99
SECTION 12.1Notices 99
SECTION 12.2Amendments, Waivers and Consents 101
SECTION 12.3Expenses; Indemnity 103
SECTION 12.4Right of Setoff 105
SECTION 12.5Governing Law; Jurisdiction, Etc 105
SECTION 12.6Waiver of Jury Trial 106
SECTION 12.7Reversal of Payments 106
SECTION 12.8Injunctive Relief 107
SECTION 12.9Successors and Assigns; Participations 107
SECTION 12.10Treatment of Certain Information; Confidentiality 110
SECTION 12.11Performance of Duties 111
SECTION 12.12All Powers Coupled with Interest 111
SECTION 12.13Survival 111
SECTION 12.14Titles and Captions 112
SECTION 12.15Severability of Provisions 112
SECTION 12.16Counterparts; Integration; Effectiveness; Electronic Execution 112
SECTION 12.17Term of Agreement 113
SECTION 12.18USA PATRIOT Act; Anti-Money Laundering Laws 113
SECTION 12.19Independent Effect of Covenants 113
SECTION 12.20No Advisory or Fiduciary Responsibility 113
SECTION 12.21Amendment and Restatement; No Novation 114
SECTION 12.22Inconsistencies with Other Documents 115
SECTION 12.23Acknowledgement and Consent to Bail-In of Affected Financial
Institutions 115
SECTION 12.24Certain ERISA Matters 115
SECTION 12.25Acknowledgement Regarding Any Supported QFCs 116
--- [START OF PAGE 30] ---
viEXHIBITS
Exhibit A-1 -Form of Revolving Credit Note
Exhibit A-2 -Form of Swingline Note
Exhibit B -Form of Notice of Borrowing
Exhibit C -Form of Notice of Account Designation
Exhibit D -Form of Notice of Prepayment
Exhibit E -Form of Notice of Conversion/Continuation
Exhibit F -Form of Officer’s Compliance Certificate
Exhibit G -Form of Assignment and Assumption
Exhibit H-1 -Form of U.S. Tax Compliance Certificate (Non-
Partnership Foreign Lenders)
Exhibit H-2 -Form of U.S. Tax Compliance Certificate (Non-
Partnership Foreign Participants)
Exhibit H-3 -Form of U.S. Tax Compliance Certificate (Foreign
Participant Partnerships)
Exhibit H-4 -Form of U.S. Tax Compliance Certificate (Foreign Lender
Partnerships)
Exhibit I -Form of Joinder Agreement
SCHEDULES
Schedule 1.1 -Commitments and Commitment Percentages
Schedule 7.1 -Jurisdictions of Organization and Qualification
Schedule 7.2 -Subsidiaries and Capitalization
Schedule 7.6 -Tax Matters
Schedule 7.9 -ERISA Plans
Schedule 7.12 -Material Contracts
Schedule 7.13 -Labor and Collective Bargaining Agreements
Schedule 7.18 -Real Property
Schedule 7.19 -Litigation
Schedule 8.20 -Post-Closing Matters
Schedule 9.1 -Existing Indebtedness
Schedule 9.2 -Existing Liens
Schedule 9.3 -Existing Loans, Advances and Investments
Schedule 9.7 -Transactions with Affiliates
Schedule 12.21 -Existing Letters of Credit
--- [START OF PAGE 31] ---
AMENDED AND RESTATED CREDIT AGREEMENT, dated as of June 28, 2023, by and
among MAYVILLE ENGINEERING COMPANY, INC., a Wisconsin corporation (“Mayville” or the
“Borrower”), the lenders who are party to this Agreement and the lenders who may become a party to
this Agreement pursuant to the terms hereof, as Lenders, and WELLS FARGO BANK, NATIONAL
ASSOCIATION, a national banking association, as Administrative Agent for the Lenders.
STATEMENT OF PURPOSE
The Borrower, the lenders party thereto, and WELLS FARGO BANK, NATIONAL
ASSOCIATION, a national banking association, as administrative agent for such lenders, are party to
that certain Amended and Restated Credit Agreement dated as of September 26, 2019 (as amended to
date, the “Existing Credit Agreement”).  
The Borrower has requested, and subject to the terms and conditions set forth in this
Agreement, the Administrative Agent and the Lenders have agreed to extend certain credit facilities
to the Borrower to be used, in part, to restructure, and to increase, the indebtedness and other
obligations outstanding under the Existing Credit Agreement and amend and restate the Existing
Credit Agreement in its entirety.
NOW, THEREFORE, for good and valuable consideration, the receipt and sufficiency of
which are hereby acknowledged by the parties hereto, such parties hereby agree as follows:
ARTICLE I
DEFINITIONS
SECTION 1.1 Definitions.  The following terms when used in this Agreement shall
have the meanings assigned to them below:
“Accu-Fab Acquisition” means the contemplated Acquisition by the Borrower of all of the
outstanding Equity Interests of Accu-Fab, LLC, a Delaware limited liability company, from Tide
Rock Yieldco, LLC, a Delaware limited liability company, pursuant to that certain Purchase
Agreement, dated as of May 23, 2025, by and among the Borrower, Accu-Fab, LLC and Tide Rock
Yieldco, LLC.  
“Acquisition” means any transaction, or any series of related transactions, consummated on
or after the date of this Agreement, by which any Credit Party or any of its 
'@
$ws.Range("D11").Value = $v_D11

$v_E11 = @'
{
  "errors": [
    {
      "location": "Page 3, Section 2.1(a)",
      "error": "The capitalized term \"Amendment\" is used but not defined. The document defines \"this Agreement\" (capitalized) to refer to the First Amendment to Amended and Restated Credit Agreement, so \"Amendment\" should likely be \"Agreement\" for consistency.",
      "suggestion": "Replace \"this Amendment\" with \"this Agreement\".",
      "exact_quote": "delivery and performance of this Amendment"
    },
    {
      "location": "Page 3, Section 3(a)(ii)",
      "error": "The sentence ends abruptly and carries over to the next page without proper continuation, indicating a drafting error or forgotten placeholder.",
      "suggestion": "Complete the sentence or ensure proper pagination and sentence flow.",
      "exact_quote": "certificates as of a recent date of the good standing of each Credit Party"
    },
    {
      "location": "Page 4, Section 3(a)(iii)",
      "error": "The sentence ends abruptly and carries over to the next page without proper continuation, indicating a drafting error or forgotten placeholder.",
      "suggestion": "Complete the sentence or ensure proper pagination and sentence flow.",
      "exact_quote": "the Borrower is in compliance with the"
    },
    {
      "location": "Page 5, Section 3(b)",
      "error": "The condition refers to 'Section 3' for representations and warranties, but the representations and warranties are actually set forth in 'Section 4'. This appears to be a numbering error.",
      "suggestion": "Change 'Section 3' to 'Section 4'.",
      "exact_quote": "The representations and warranties set forth in Section 3 shall be true and correct."
    },
    {
      "location": "Page 9",
      "error": "A placeholder indicates that signature pages are expected to follow, which is a drafting error in a completed document.",
      "suggestion": "Remove the placeholder or replace it with appropriate text if the document is finalized.",
      "exact_quote": "[Signature pages to follow]"
    },
    {
      "location": "Page 47, Article I (Definitions)",
      "error": "There is a pagination error. The page number is '13', but it should sequentially be '14' after the preceding page (Page 46) which was '13'.",
      "suggestion": "Correct the page number to '14'.",
      "exact_quote": "13"
    },
    {
      "location": "Page 49, Article I (Definitions)",
      "error": "There is a pagination error. The page number is '14', but it should sequentially be '15' after the preceding page (Page 48) which was '14'.",
      "suggestion": "Correct the page number to '15'.",
      "exact_quote": "14"
    },
    {
      "location": "Page 54, Article I (Definitions)",
      "error": "The singular term 'Note' is capitalized and used in the definition of 'Loan Documents' but is not explicitly defined. While 'Notes' (plural) is defined later, consistency requires either 'Note' to be defined or 'Notes' to be used consistently.",
      "suggestion": "Define 'Note' as the singular of 'Notes' or replace 'each Note' with 'the Notes' in the definition of 'Loan Documents'.",
      "exact_quote": "each Note"
    },
    {
      "location": "Page 55, Article I (Definitions)",
      "error": "The capitalized term 'Outstanding Amount' is used in the definition of 'Minimum Collateral Amount' but is not defined.",
      "suggestion": "Add a definition for 'Outstanding Amount'.",
      "exact_quote": "Outstanding Amount"
    },
    {
      "location": "Page 60, Article I (Definitions)",
      "error": "There is a typographical error in 'Sanctioned Peron(s)'. It should be 'Sanctioned Person(s)'.",
      "suggestion": "Correct 'Peron(s)' to 'Person(s)'.",
      "exact_quote": "Sanctioned Peron(s)"
    },
    {
      "location": "Page 63, Article I (Definitions)",
      "error": "The phrase 'for the ratable benefit and the Secured Parties' contains a grammatical error. It should likely be 'for the ratable benefit of the Secured Parties'.",
      "suggestion": "Change 'benefit and the Secured Parties' to 'benefit of the Secured Parties'.",
      "exact_quote": "for the ratable benefit and the Secured Parties"
    },
    {
      "location": "Page 65, Article I (Definitions)",
      "error": "There is a pagination error. The page number is '29', but it should sequentially be '30' after the preceding page (Page 63) which was '28'.",
      "suggestion": "Correct the page number to '30'.",
      "exact_quote": "29"
    },
    {
      "location": "Page 65, Article I (Definitions)",
      "error": "The term 'Base Rate SOFR Determination Day' is used, but the defined term is 'Base Rate Term SOFR Determination Day'. This is an inconsistency or typo.",
      "suggestion": "Replace 'Base Rate SOFR Determination Day' with 'Base Rate Term SOFR Determination Day'.",
      "exact_quote": "Base Rate SOFR Determination Day"
    },
    {
      "location": "Page 66, Article I (Definitions)",
      "error": "There is a pagination error. The page number is '29', but it should sequentially be '31' after the preceding page (Page 65) which was '29'.",
      "suggestion": "Correct the page number to '31'.",
      "exact_quote": "29"
    },
    {
      "location": "Page 67, Article I (Definitions)",
      "error": "Section 4.2(a) is referenced in the definition of 'U.S. Government Securities Business Day', but Article IV is explicitly 'RESERVED' in this document and does not contain a Section 4.2(a). This is likely an incorrect cross-reference.",
      "suggestion": "Correct the cross-reference to the appropriate section, e.g., Section 6.2(a).",
      "exact_quote": "Sections 2.3(a), 2.4(c), 4.2(a) and 5.2"
    },
    {
      "location": "Page 75, Section 2.3(a)",
      "error": "The text uses '(1) U.S. Government Securities Business Day' where a numerical word 'one' is expected, consistent with other numerical spellings in parentheses (e.g., 'three (3)').",
      "suggestion": "Change '(1)' to 'one'.",
      "exact_quote": "no later than (1) U.S. Government Securities Business Day"
    },
    {
      "location": "Page 77, Section 2.4(a) and (b)",
      "error": "Paragraph (b) 'Mandatory Prepayments' is incorrectly formatted as a continuation of the sentence from paragraph (a) 'Repayment on Termination Date' instead of as a new, distinct sub-section.",
      "suggestion": "Start ' (b) Mandatory Prepayments.' on a new line and ensure it is properly formatted as a separate sub-section.",
      "exact_quote": "accrued but unpaid interest thereon (b) Mandatory Prepayments."
    },
    {
      "location": "Page 77, Section 2.4(d)",
      "error": "The section is marked 'Reserved', indicating a placeholder or an incomplete section.",
      "suggestion": "Either provide content for this section or explicitly state that it is intentionally left blank if it is not a drafting oversight.",
      "exact_quote": "(d) Reserved."
    },
    {
      "location": "Page 92, Section 5.9",
      "error": "Section 4.4(a) is referenced in the indemnity clause, but Article IV is explicitly 'RESERVED' in this document and does not contain a Section 4.4(a). This is likely an incorrect cross-reference.",
      "suggestion": "Correct the cross-reference to the appropriate section, e.g., Section 2.4(c).",
      "exact_quote": "Section 4.4(a)"
    }
  ]
}
'@
$ws.Range("E11").Value = $v_E11

$v_A12 = @'
2025-12-30 01:28:49
'@
$ws.Range("A12").Value = $v_A12

$v_B12 = @'
gemini-2.5-flash
'@
$ws.Range("B12").Value = $v_B12

$v_C12 = @'
d6568581-5567-430e-a7af-ad98d7fd75de.pdf
'@
$ws.Range("C12").Value = $v_C12

$v_D12 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document fragment.
Input: The attached text from a contract.
CRITICAL INSTRUCTIONS:
1. **Assume Isolation with Common Sense**: Do NOT assume missing definitions exist in other documents. However, IGNORE common commercial lending terms typically defined in a base Credit Agreement (e.g., "Borrower", "Administrative Agent", "Lender", "Business Day", "Dollars", "GAAP", "Material Adverse Effect"). Only flag specific, deal-specific, or unusual capitalized terms that are undefined.
2. **Logic Check:** Check all math and logic tables.
3. **Drafting Errors:** Find any placeholders like "[__]" or blank lines that were forgotten.
4. **Exact Quotes:** For every error, you MUST provide the `exact_quote` from the text that contains the error. This is used for highlighting.
Output Format:
Return ONLY a valid JSON object with the following structure:
{
  "errors": [
    {
      "location": "Page 3, Section 2.1",
      "error": "Description of the error",
      "suggestion": "Suggested fix",
      "exact_quote": "Exact text substring to highlight in red"
    }
  ]
}
If no errors are found, return {"errors": []}.
--- CONTRACT TEXT BEGINS ---
--- [START OF PAGE 1] ---
 EX-4.1 3 d45101 1dex41.htm EX-4.1 
 Exhibit 4.1 
 CAPIT AL ONE PRIME AUT O RECEIV ABLES TRUST 2023-1 
 Class A-1 4.900% Auto Loan Asset Backed Notes 
 Class A-2 5.20% Auto Loan Asset Backed Notes 
 Class A-3 4.87% Auto Loan Asset Backed Notes 
 Class A-4 4.76% Auto Loan Asset Backed Notes 
 Class B 5.07% Auto Loan Asset Backed Notes 
 Class C 5.36% Auto Loan Asset Backed Notes 
 Class D 6.34% Auto Loan Asset Backed Notes 
 INDENTURE 
 Dated as of February 23, 2023 
 WILMINGT ON TRUST , NATIONAL ASSOCIA TION, 
 as the Indenture T rustee 
--- [START OF PAGE 2] ---
 CROSS REFERENCE T ABLE  1 
 TIA  Section  Indenture  Section 
 310  (a) (1)  6.11 
 (a) (2)  6.11 
 (a) (3)  6.10; 6.1 1 
 (a) (4)  N.A.  2 
 (a) (5)  6.11 
 (b)  6.8; 6.1 1 
 (c)  N.A. 
 311  (a)  6.12 
 (b)  6.12 
 (c)  N.A. 
 312  (a)  7.1 
 (b)  7.2 
 (c)  7.2 
 313  (a)  7.3 
 (b) (1)  7.3 
 (b) (2)  7.3 
 (c)  7.3 
 (d)  7.3 
 314  (a)  3.9 
 (b)  11.1, 3.6 
 (c) (1)  11.1 
 (c) (2)  11.1 
 (c) (3)  11.1 
 (d)  11.1 
 (e)  11.1 
 (f)  N.A. 
 315  (a)  6.1(b) 
 (b)  6.5 
 (c)  6.1(a) 
 (d)  6.1(c) 
 (e)  5.13 
 316  (a) (1) (A)  5.11 
 (a) (1) (B)  5.12 
 (a) (2)  N.A. 
 (b)  5.7 
 (c)  5.6(b) 
 317  (a) (1)  5.3(b) 
 (a) (2)  5.3(d) 
 (b)  3.3 
 318  (a)  11.7 
 1  Note: This Cross Reference T able shall not, for  any purpose, be deemed to be part of this Indenture. 
 2  N.A. means Not Applicable. 
 COPAR 2023-1 Indentur e 
--- [START OF PAGE 3] ---
 Page 
 ARTICLE I  DEFINITIONS AND INCORPORA TION BY REFERENCE  2 
 Section 1.1  Definitions  2 
 Section 1.2  Incorporation by Reference of T rust Indenture Act  2 
 Section 1.3  Other Interpretive Provisions  2 
 ARTICLE II  THE NOTES  3 
 Section 2.1  Form  3 
 Section 2.2  Execution, Authentication and Delivery  3 
 Section 2.3  Temporary Notes  4 
 Section 2.4  Registration of T ransfer and Exchange  4 
 Section 2.5  Mutilated, Destroyed, Lost or Stolen Notes  6 
 Section 2.6  Persons Deemed Owners  6 
 Section 2.7  Payment of Principal and Interest; Defaulted Interest  7 
 Section 2.8  Cancellation  7 
 Section 2.9  Release of Collateral  8 
 Section 2.10  Book-Entry Notes  8 
 Section 2.1 1  Notices to Clearing Agency  9 
 Section 2.12  Definitive Notes  9 
 Section 2.13  Authenticating Agents  10 
 Section 2.14  Paying Agent  10 
 Section 2.15  Tax and Accounting T reatment  11 
 Section 2.16  Certain T ransfer Restrictions on all Classes of the Notes  11 
 Section 2.17  Certain T ransfer Restrictions  12 
 Section 2.18  Transfer Restrictions on Certain Notes Upon a Sale of a Certificate  12 
 Section 2.19  Certain T ransfer Restrictions on the 144A Notes  13 
 ARTICLE III  COVENANTS  17 
 Section 3.1  Payment of Principal and Interest  17 
 Section 3.2  Maintenance of Of fice or Agency  17 
 Section 3.3  Money for Payments to Be Held in T rust  17 
 Section 3.4  Existence  19 
 Section 3.5  Protection of Collateral  19 
--- [START OF PAGE 4] ---
--- [START OF PAGE 5] ---
 Page 
 Section 3.6  Opinions as to Collateral  20 
 Section 3.7  Performance of Obligations; Servicing of Receivables  20 
 Section 3.8  Negative Covenants  21 
 Section 3.9  Annual Compliance Statement  22 
 Section 3.10  Restrictions on Certain Other Activities  23 
 Section 3.1 1  Restricted Payments  23 
 Section 3.12  Notice of Events of Default  23 
 Section 3.13  Further Instruments and Acts  23 
 Section 3.14  Compliance with Laws  23 
 Section 3.15  Removal of Administrator  23 
 Section 3.16  Perfection Representations, W arranties and Covenants  24 
 Section 3.17  Investment Company Act Representation  24 
 ARTICLE IV  SATISF ACTION AND DISCHARGE  24 
 Section 4.1  Satisfaction and Dischar ge of Indenture  24 
 Section 4.2  Application of T rust Money  25 
 Section 4.3  Repayment of Monies Held by Paying Agent  25 
 ARTICLE V  REMEDIES  25 
 Section 5.1  Events of Default  25 
 Section 5.2  Acceleration of Maturity; W aiver of Event of Default  26 
 Section 5.3  Collection of Indebtedness and Suits for Enforcement by the Indenture T rustee  27 
 Section 5.4  Remedies; Priorities  29 
 Section 5.5  Optional Preservation of the Collateral  32 
 Section 5.6  Limitation of Suits  32 
 Section 5.7  Rights of Noteholders to Receive Principal and Interest  33 
 Section 5.8  Restoration of Rights and Remedies  33 
 Section 5.9  Rights and Remedies Cumulative  33 
 Section 5.10  Delay or Omission Not a W aiver  33 
 Section 5.1 1  Control by Noteholders  34 
 Section 5.12  Waiver of Past Defaults  34 
 Section 5.13  Undertaking for Costs  35 
--- [START OF PAGE 6] ---
--- [START OF PAGE 7] ---
 Page 
 Section 5.14  Waiver of Stay or Extension Laws  35 
 Section 5.15  Action on Notes  35 
 Section 5.16  Performance and Enforcement of Certain Obligations  35 
 Section 5.17  Sale of Collateral  36 
 ARTICLE VI  THE INDENTURE TRUSTEE  37 
 Section 6.1  Duties of the Indenture T rustee  37 
 Section 6.2  Rights of the Indenture T rustee  38 
 Section 6.3  Individual Rights of the Indenture T rustee  40 
 Section 6.4  The Indenture T rustee’ s Disclaimer  40 
 Section 6.5  Notice of Defaults  40 
 Section 6.6  Reports by the Paying Agent  40 
 Section 6.7  Compensation and Indemnity  40 
 Section 6.8  Removal, Resignation and Replacement of the Indenture T rustee  41 
 Section 6.9  Successor Indenture T rustee by Mer ger  42 
 Section 6.10  Appointment of Co-Indenture T rustee or Separate Indenture T rustee  43 
 Section 6.1 1  Eligibility; Disqualification  44 
 Section 6.12  Preferential Collection of Claims Against the Issuer  44 
 Section 6.13  Representations and W arranties  44 
 ARTICLE VII  NOTEHOLDERS’ LISTS AND REPOR TS  45 
 Section 7.1  The Issuer to Furnish the Indenture T rustee Names and Addresses of Noteholders  45 
 Section 7.2  Preservation of Information; Communications to Noteholders  45 
 Section 7.3  Reports by the Indenture T rustee  45 
 Section 7.4  Statements to Certificateholders and Noteholders  45 
 Section 7.5  Noteholder Demand for Repurchase, Dispute Resolution  47 
 Section 7.6  Investor Action to Initiate an Asset Review  48 
 ARTICLE VIII  ACCOUNTS, DISBURSEMENTS AND RELEASES  49 
 Section 8.1  Collection of Money  49 
 Section 8.2  Trust Accounts  49 
--- [START OF PAGE 8] ---
 Page 
 Section 8.3  General Provisions Regarding Accounts  51 
 Section 8.4  Additional W ithdrawals and Deposits  54 
 Section 8.5  Distributions  54 
 Section 8.6  Release of Collateral  56 
 Section 8.7  Opinion of Counsel  56 
 ARTICLE IX  SUPPLEMENT AL INDENTURES  57 
 Section 9.1  Supplemental Indentures W ithout Consent of Noteholders  57 
 Section 9.2  Supplemental Indentures with Consent of Noteholders  58 
 Section 9.3  Execution of Supplemental Indentures  59 
 Section 9.4  Effect of Supplemental Indenture  60 
 Section 9.5  Conformity with T rust Indenture Act  60 
 Section 9.6  Reference in Notes to Supplemental Indentures  60 
 ARTICLE X  REDEMPTION OF NOTES  60 
 Section 10.1  Redemption  60 
 Section 10.2  Form of Redemption Notice  61 
 Section 10.3  Notes Payable on Redemption Date  61 
 ARTICLE XI  MISCELLANEOUS  62 
 Section 1 1.1  Compliance Certificates and Opinions, etc.  62 
 Section 1 1.2  Form of Documents Delivered to the Indenture T rustee  63 
 Section 1 1.3  Acts of Noteholders  64 
 Section 1 1.4  Notices  65 
 Section 1 1.5  Notices to Noteholders; W aiver  65 
 Section 1 1.6  Alternate Payment and Notice Provisions  65 
 Section 1 1.7  Conflict with T rust Indenture Act  65 
 Section 1 1.8  Information Requests  66 
 Section 1 1.9  Effect of Headings and T able of Contents  66 
 Section 1 1.10  Successors and Assigns  66 
 Section 1 1.11  Separability  66 
 Section 1 1.12  Benefits of Indenture  66 
 Section 1 1.13  Legal Holidays  66 
--- [START OF PAGE 9] ---
--- [START OF PAGE 10] ---
 Page 
 Section 1 1.14  GOVERNING LA W  66 
 Section 1 1.15  Counterparts  66 
 Section 1 1.16  Recording of Indenture  67 
 Section 1 1.17  Trust Obligation  67 
 Section 1 1.18  No Petition  67 
 Section 1 1.19  Submission to Jurisdiction; W aiver of Jury T rial  68 
 Section 1 1.20  Subordination of Claims  68 
 Section 1 1.21  U.S.A Patriot Act  69 
 Section 1 1.22  Beneficial Ownership  69 
 Section 1 1.23  Limitation of Liability  69 
 ARTICLE XII  COMPLIANCE WITH THE FDIC RULE  70 
 Section 12.1  Purpose  70 
 Section 12.2  Requirements of the FDIC Rule  71 
 Section 12.3  Performance  72 
 Section 12.4  Actions Upon Repudiation  73 
 Section 12.5  Notice  75 
 Section 12.6  Reservation of Rights  75 
 Schedule I  Perfection Representations, W arranties and Covenants 
 Exhibit A  Form of Notes 
--- [START OF PAGE 11] ---
 This INDENTURE, dated as of February 28, 2023 (as amended, supplemented, or otherwise modified and in ef fect from time to time, this 
 “  Indenture  ”), is between CAPIT AL ONE PRIME AUT O RECEIV ABLES  TRUST 2023-1, a Delaware statutory trust (the “  Issuer  ”),  and W ilmington 
 Trust, National Association, a national banking association, solely as trustee and not in its individual capacity (the “  Indenture T rustee  ”). 
 Each party agrees as follows for the benefit of the other party and the equal and ratable benefit of the Holders of the Issuer ’s Class A-1 
 4.900% Auto Loan Asset Backed Notes (the “  Class A-1  Notes  ”), Class A-2 5.20% Auto Loan Asset Backed Notes  (the “  Class A-2 Notes  ”), Class 
 A-3 4.87% Auto Loan Asset Backed Notes (the “  Class  A-3 Notes  ”), Class A-4 4.76% Auto Loan Asset Backed  Notes (the “  Class A-4 Notes  ”, and 
 together with the Class A-1 Notes, the Class A-2 Notes and the Class A-3 Notes, the “  Class A Notes  ”), Class  B 5.07% Auto Loan Asset Backed 
 Notes (the “  Class B Notes  ”), Class C 5.36% Auto Loan  Asset Backed Notes (the “  Class C Notes  ”) and Class  D 6.34% Auto Loan Asset Backed 
 Notes (the “  Class D Notes  ”, and together with the  Class A-1 Notes, the Class A-2 Notes, the Class A-3 Notes, the Class A-4 Notes, the Class B Notes 
 and the Class C Notes, the “  Notes  ”). 
 GRANTING CLAUSE 
 The Issuer , to secure the payment of principal of and interest on, and any other amounts owing in respect of, the Notes equally and ratably 
 without prejudice, priority or distinction except as set forth herein, and to secure compliance with the provisions of this Indenture, hereby Grants in trust 
 to the Indenture T rustee on the Closing Date, as trustee for the benefit of the Noteholders, all of the Issuer ’s right, title and interest, whether now owned 
 or hereafter acquired, in and to (i) the T rust Estate and (ii) all present and future claims, demands, causes and choses in action in respect of any or all of 
 the T rust Estate and all payments on or under and all proceeds of every kind and nature whatsoever in respect of any or all of the T rust Estate, including 
 all proceeds of the conversion, voluntary or involuntary , into cash or other liquid property , all cash proceeds, accounts, accounts receivable, notes, drafts, 
 acceptances, chattel paper , checks, deposit accounts, insurance proceeds, condemnation awards, rights to payment of any and every kind and other forms 
 of obligations and receivables, instruments, securities, financial assets and other property which at any time constitute all or part of or are included in the 
 proceeds of any of the T rust Estate (collectively , the “  Collateral  ”). 
 The Indenture T rustee, on behalf of the Noteholders, acknowledges the foregoing Grant, accepts the trusts under this Indenture and 
 agrees to perform its duties required in this Indenture in accordance with the provisions of this Indenture. 
 The foregoing Grant is made in trust to secure (i) the payment of principal of and interest on, and any other amounts owing in respect of, 
 the Notes, equally and ratably without prejudice, priority or distinction except as set forth herein and (ii) compliance with the provisions of this 
 Indenture, all as provided in this Indenture. 
 Without limiting the foregoing Grant, any Receivable repurchased or purchased by (a) the Servicer pursuant to  Section 3.6  of the Servicing 
 Agreement or (b) by the Bank pursuant to  Section 3.4  of the Purchase Agreement shall be deemed to be automatically released from the lien of this 
 Indenture without any action being taken by the Indenture T rustee upon payment by the applicable purchaser of the related Repurchase Price for such 
 Repurchased Receivable. 
 COPAR 2023-1 Indentur e 
--- [START OF PAGE 12] ---
 ARTICLE I 
 DEFINITIONS AND INCORPORA TION BY REFERENCE 
 SECTION 1.1  Definitions  . Capitalized terms are used  in this Indenture as defined in  A  p  p  endix A  to the  Sale Agreement, dated as of the 
 date hereof (as amended, supplemented, or otherwise modified and in ef fect from time to time, the “  Sale  Agreement  ”), between the Issuer and Capital 
 One Auto Receivables, LLC, which also contains rules as to usage that are applicable herein. 
 SECTION 1.2  Incorporation by Reference of T rust Indenture  Act  . Whenever this Indenture refers to a provision  of the TIA, the provision 
 is incorporated by reference in and made a part of this Indenture. The following TIA terms used in this Indenture have the following meanings: 
 “  Commission  ” means the Securities and Exchange Commission. 
 “  indenture securities  ” means the Notes. 
 “  indenture securi  ty  holder  ” means a Noteholder . 
 “  indenture to be qualified  ” means this Indenture. 
 “  indenture trustee  ” or “  institutional trustee  ” means  the Indenture T rustee. 
 “  obl  ig  or  ” on the indenture securities means the Issuer  and any other obligor on the indenture securities. 
 All other TIA terms used in this Indenture that are defined by the TIA, defined by TIA reference to another statute or defined by 
 Commission rule have the meaning assigned to them by such definitions. 
 SECTION 1.3  Other Interpretive Provisions  . All terms  defined in this Indenture shall have the defined meanings when used in any 
 certificate or other document delivered pursuant hereto unless otherwise defined therein. For purposes of this Indenture and all such certificates and 
 other documents, unless the context otherwise requires: (a) accounting terms not otherwise defined in this Indenture, and accounting terms partly 
 defined in this Indenture to the extent not defined, shall have the respective meanings given to them under GAAP (  provided  , that, to the extent that the 
 definitions in this Indenture and GAAP conflict, the definitions in this Indenture shall control); (b) terms defined in Article 9 of the UCC as in ef fect in 
 the relevant jurisdiction and not otherwise defined in this Indenture are used as defined in that Article; (c) the words “hereof,” “herein” and “hereunder” 
 and words of similar import refer to this Indenture as a whole and not to any particular provision of this Indenture; (d) references to any Article, 
 Section, Schedule, Appendix or Exhibit are references to Articles, Sections, Schedules, Appendices and Exhibits in or to this Indenture and references 
 to any paragraph, subsection, clause or other subdivision within any Section or definition refer to such paragraph, subsection, clause or other 
--- CONTRACT TEXT ENDS ---
'@
$ws.Range("D12").Value = $v_D12

$v_E12 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document fragment.
Input: The attached text from a contract.
CRITICAL INSTRUCTIONS:
1. **Assume Isolation with Common Sense**: Do NOT assume missing definitions exist in other documents. However, IGNORE common commercial lending terms typically defined in a base Credit Agreement (e.g., "Borrower", "Administrative Agent", "Lender", "Business Day", "Dollars", "GAAP", "Material Adverse Effect"). Only flag specific, deal-specific, or unusual capitalized terms that are undefined.
2. **Logic Check:** Check all math and logic tables.
3. **Drafting Errors:** Find any placeholders like "[__]" or blank lines that were forgotten.
4. **Exact Quotes:** For every error, you MUST provide the `exact_quote` from the text that contains the error. This is used for highlighting.
Output Format:
Return ONLY a valid JSON object with the following structure:
{
  "errors": [
    {
      "location": "Page 3, Section 2.1",
      "error": "Description of the error",
      "suggestion": "Suggested fix",
      "exact_quote": "Exact text substring to highlight in red"
    }
  ]
}
If no errors are found, return {"errors": []}.
--- CONTRACT TEXT BEGINS ---
--- [START OF PAGE 1] ---
 EX-4.1 3 d45101 1dex41.htm EX-4.1 
 Exhibit 4.1 
 CAPIT AL ONE PRIME AUT O RECEIV ABLES TRUST 2023-1 
 Class A-1 4.900% Auto Loan Asset Backed Notes 
 Class A-2 5.20% Auto Loan Asset Backed Notes 
 Class A-3 4.87% Auto Loan Asset Backed Notes 
 Class A-4 4.76% Auto Loan Asset Backed Notes 
 Class B 5.07% Auto Loan Asset Backed Notes 
 Class C 5.36% Auto Loan Asset Backed Notes 
 Class D 6.34% Auto Loan Asset Backed Notes 
 INDENTURE 
 Dated as of February 23, 2023 
 WILMINGT ON TRUST , NATIONAL ASSOCIA TION, 
 as the Indenture T rustee 
--- [START OF PAGE 2] ---
 CROSS REFERENCE T ABLE  1 
 TIA  Section  Indenture  Section 
 310  (a) (1)  6.11 
 (a) (2)  6.11 
 (a) (3)  6.10; 6.1 1 
 (a) (4)  N.A.  2 
 (a) (5)  6.11 
 (b)  6.8; 6.1 1 
 (c)  N.A. 
 311  (a)  6.12 
 (b)  6.12 
 (c)  N.A. 
 312  (a)  7.1 
 (b)  7.2 
 (c)  7.2 
 313  (a)  7.3 
 (b) (1)  7.3 
 (b) (2)  7.3 
 (c)  7.3 
 (d)  7.3 
 314  (a)  3.9 
 (b)  11.1, 3.6 
 (c) (1)  11.1 
 (c) (2)  11.1 
 (c) (3)  11.1 
 (d)  11.1 
 (e)  11.1 
 (f)  N.A. 
 315  (a)  6.1(b) 
 (b)  6.5 
 (c)  6.1(a) 
 (d)  6.1(c) 
 (e)  5.13 
 316  (a) (1) (A)  5.11 
 (a) (1) (B)  5.12 
 (a) (2)  N.A. 
 (b)  5.7 
 (c)  5.6(b) 
 317  (a) (1)  5.3(b) 
 (a) (2)  5.3(d) 
 (b)  3.3 
 318  (a)  11.7 
 1  Note: This Cross Reference T able shall not, for  any purpose, be deemed to be part of this Indenture. 
 2  N.A. means Not Applicable. 
 COPAR 2023-1 Indentur e 
--- [START OF PAGE 3] ---
 Page 
 ARTICLE I  DEFINITIONS AND INCORPORA TION BY REFERENCE  2 
 Section 1.1  Definitions  2 
 Section 1.2  Incorporation by Reference of T rust Indenture Act  2 
 Section 1.3  Other Interpretive Provisions  2 
 ARTICLE II  THE NOTES  3 
 Section 2.1  Form  3 
 Section 2.2  Execution, Authentication and Delivery  3 
 Section 2.3  Temporary Notes  4 
 Section 2.4  Registration of T ransfer and Exchange  4 
 Section 2.5  Mutilated, Destroyed, Lost or Stolen Notes  6 
 Section 2.6  Persons Deemed Owners  6 
 Section 2.7  Payment of Principal and Interest; Defaulted Interest  7 
 Section 2.8  Cancellation  7 
 Section 2.9  Release of Collateral  8 
 Section 2.10  Book-Entry Notes  8 
 Section 2.1 1  Notices to Clearing Agency  9 
 Section 2.12  Definitive Notes  9 
 Section 2.13  Authenticating Agents  10 
 Section 2.14  Paying Agent  10 
 Section 2.15  Tax and Accounting T reatment  11 
 Section 2.16  Certain T ransfer Restrictions on all Classes of the Notes  11 
 Section 2.17  Certain T ransfer Restrictions  12 
 Section 2.18  Transfer Restrictions on Certain Notes Upon a Sale of a Certificate  12 
 Section 2.19  Certain T ransfer Restrictions on the 144A Notes  13 
 ARTICLE III  COVENANTS  17 
 Section 3.1  Payment of Principal and Interest  17 
 Section 3.2  Maintenance of Of fice or Agency  17 
 Section 3.3  Money for Payments to Be Held in T rust  17 
 Section 3.4  Existence  19 
 Section 3.5  Protection of Collateral  19 
--- [START OF PAGE 4] ---
--- [START OF PAGE 5] ---
 Page 
 Section 3.6  Opinions as to Collateral  20 
 Section 3.7  Performance of Obligations; Servicing of Receivables  20 
 Section 3.8  Negative Covenants  21 
 Section 3.9  Annual Compliance Statement  22 
 Section 3.10  Restrictions on Certain Other Activities  23 
 Section 3.1 1  Restricted Payments  23 
 Section 3.12  Notice of Events of Default  23 
 Section 3.13  Further Instruments and Acts  23 
 Section 3.14  Compliance with Laws  23 
 Section 3.15  Removal of Administrator  23 
 Section 3.16  Perfection Representations, W arranties and Covenants  24 
 Section 3.17  Investment Company Act Representation  24 
 ARTICLE IV  SATISF ACTION AND DISCHARGE  24 
 Section 4.1  Satisfaction and Dischar ge of Indenture  24 
 Section 4.2  Application of T rust Money  25 
 Section 4.3  Repayment of Monies Held by Paying Agent  25 
 ARTICLE V  REMEDIES  25 
 Section 5.1  Events of Default  25 
 Section 5.2  Acceleration of Maturity; W aiver of Event of Default  26 
 Section 5.3  Collection of Indebtedness and Suits for Enforcement by the Indenture T rustee  27 
 Section 5.4  Remedies; Priorities  29 
 Section 5.5  Optional Preservation of the Collateral  32 
 Section 5.6  Limitation of Suits  32 
 Section 5.7  Rights of Noteholders to Receive Principal and Interest  33 
 Section 5.8  Restoration of Rights and Remedies  33 
 Section 5.9  Rights and Remedies Cumulative  33 
 Section 5.10  Delay or Omission Not a W aiver  33 
 Section 5.1 1  Control by Noteholders  34 
 Section 5.12  Waiver of Past Defaults  34 
 Section 5.13  Undertaking for Costs  35 
--- [START OF PAGE 6] ---
--- [START OF PAGE 7] ---
 Page 
 Section 5.14  Waiver of Stay or Extension Laws  35 
 Section 5.15  Action on Notes  35 
 Section 5.16  Performance and Enforcement of Certain Obligations  35 
 Section 5.17  Sale of Collateral  36 
 ARTICLE VI  THE INDENTURE TRUSTEE  37 
 Section 6.1  Duties of the Indenture T rustee  37 
 Section 6.2  Rights of the Indenture T rustee  38 
 Section 6.3  Individual Rights of the Indenture T rustee  40 
 Section 6.4  The Indenture T rustee’ s Disclaimer  40 
 Section 6.5  Notice of Defaults  40 
 Section 6.6  Reports by the Paying Agent  40 
 Section 6.7  Compensation and Indemnity  40 
 Section 6.8  Removal, Resignation and Replacement of the Indenture T rustee  41 
 Section 6.9  Successor Indenture T rustee by Mer ger  42 
 Section 6.10  Appointment of Co-Indenture T rustee or Separate Indenture T rustee  43 
 Section 6.1 1  Eligibility; Disqualification  44 
 Section 6.12  Preferential Collection of Claims Against the Issuer  44 
 Section 6.13  Representations and W arranties  44 
 ARTICLE VII  NOTEHOLDERS’ LISTS AND REPOR TS  45 
 Section 7.1  The Issuer to Furnish the Indenture T rustee Names and Addresses of Noteholders  45 
 Section 7.2  Preservation of Information; Communications to Noteholders  45 
 Section 7.3  Reports by the Indenture T rustee  45 
 Section 7.4  Statements to Certificateholders and Noteholders  45 
 Section 7.5  Noteholder Demand for Repurchase, Dispute Resolution  47 
 Section 7.6  Investor Action to Initiate an Asset Review  48 
 ARTICLE VIII  ACCOUNTS, DISBURSEMENTS AND RELEASES  49 
 Section 8.1  Collection of Money  49 
 Section 8.2  Trust Accounts  49 
--- [START OF PAGE 8] ---
 Page 
 Section 8.3  General Provisions Regarding Accounts  51 
 Section 8.4  Additional W ithdrawals and Deposits  54 
 Section 8.5  Distributions  54 
 Section 8.6  Release of Collateral  56 
 Section 8.7  Opinion of Counsel  56 
 ARTICLE IX  SUPPLEMENT AL INDENTURES  57 
 Section 9.1  Supplemental Indentures W ithout Consent of Noteholders  57 
 Section 9.2  Supplemental Indentures with Consent of Noteholders  58 
 Section 9.3  Execution of Supplemental Indentures  59 
 Section 9.4  Effect of Supplemental Indenture  60 
 Section 9.5  Conformity with T rust Indenture Act  60 
 Section 9.6  Reference in Notes to Supplemental Indentures  60 
 ARTICLE X  REDEMPTION OF NOTES  60 
 Section 10.1  Redemption  60 
 Section 10.2  Form of Redemption Notice  61 
 Section 10.3  Notes Payable on Redemption Date  61 
 ARTICLE XI  MISCELLANEOUS  62 
 Section 1 1.1  Compliance Certificates and Opinions, etc.  62 
 Section 1 1.2  Form of Documents Delivered to the Indenture T rustee  63 
 Section 1 1.3  Acts of Noteholders  64 
 Section 1 1.4  Notices  65 
 Section 1 1.5  Notices to Noteholders; W aiver  65 
 Section 1 1.6  Alternate Payment and Notice Provisions  65 
 Section 1 1.7  Conflict with T rust Indenture Act  65 
 Section 1 1.8  Information Requests  66 
 Section 1 1.9  Effect of Headings and T able of Contents  66 
 Section 1 1.10  Successors and Assigns  66 
 Section 1 1.11  Separability  66 
 Section 1 1.12  Benefits of Indenture  66 
 Section 1 1.13  Legal Holidays  66 
--- [START OF PAGE 9] ---
--- [START OF PAGE 10] ---
 Page 
 Section 1 1.14  GOVERNING LA W  66 
 Section 1 1.15  Counterparts  66 
 Section 1 1.16  Recording of Indenture  67 
 Section 1 1.17  Trust Obligation  67 
 Section 1 1.18  No Petition  67 
 Section 1 1.19  Submission to Jurisdiction; W aiver of Jury T rial  68 
 Section 1 1.20  Subordination of Claims  68 
 Section 1 1.21  U.S.A Patriot Act  69 
 Section 1 1.22  Beneficial Ownership  69 
 Section 1 1.23  Limitation of Liability  69 
 ARTICLE XII  COMPLIANCE WITH THE FDIC RULE  70 
 Section 12.1  Purpose  70 
 Section 12.2  Requirements of the FDIC Rule  71 
 Section 12.3  Performance  72 
 Section 12.4  Actions Upon Repudiation  73 
 Section 12.5  Notice  75 
 Section 12.6  Reservation of Rights  75 
 Schedule I  Perfection Representations, W arranties and Covenants 
 Exhibit A  Form of Notes 
--- [START OF PAGE 11] ---
 This INDENTURE, dated as of February 28, 2023 (as amended, supplemented, or otherwise modified and in ef fect from time to time, this 
 “  Indenture  ”), is between CAPIT AL ONE PRIME AUT O RECEIV ABLES  TRUST 2023-1, a Delaware statutory trust (the “  Issuer  ”),  and W ilmington 
 Trust, National Association, a national banking association, solely as trustee and not in its individual capacity (the “  Indenture T rustee  ”). 
 Each party agrees as follows for the benefit of the other party and the equal and ratable benefit of the Holders of the Issuer ’s Class A-1 
 4.900% Auto Loan Asset Backed Notes (the “  Class A-1  Notes  ”), Class A-2 5.20% Auto Loan Asset Backed Notes  (the “  Class A-2 Notes  ”), Class 
 A-3 4.87% Auto Loan Asset Backed Notes (the “  Class  A-3 Notes  ”), Class A-4 4.76% Auto Loan Asset Backed  Notes (the “  Class A-4 Notes  ”, and 
 together with the Class A-1 Notes, the Class A-2 Notes and the Class A-3 Notes, the “  Class A Notes  ”), Class  B 5.07% Auto Loan Asset Backed 
 Notes (the “  Class B Notes  ”), Class C 5.36% Auto Loan  Asset Backed Notes (the “  Class C Notes  ”) and Class  D 6.34% Auto Loan Asset Backed 
 Notes (the “  Class D Notes  ”, and together with the  Class A-1 Notes, the Class A-2 Notes, the Class A-3 Notes, the Class A-4 Notes, the Class B Notes 
 and the Class C Notes, the “  Notes  ”). 
 GRANTING CLAUSE 
 The Issuer , to secure the payment of principal of and interest on, and any other amounts owing in respect of, the Notes equally and ratably 
 without prejudice, priority or distinction except as set forth herein, and to secure compliance with the provisions of this Indenture, hereby Grants in trust 
 to the Indenture T rustee on the Closing Date, as trustee for the benefit of the Noteholders, all of the Issuer ’s right, title and interest, whether now owned 
 or hereafter acquired, in and to (i) the T rust Estate and (ii) all present and future claims, demands, causes and choses in action in respect of any or all of 
 the T rust Estate and all payments on or under and all proceeds of every kind and nature whatsoever in respect of any or all of the T rust Estate, including 
 all proceeds of the conversion, voluntary or involuntary , into cash or other liquid property , all cash proceeds, accounts, accounts receivable, notes, drafts, 
 acceptances, chattel paper , checks, deposit accounts, insurance proceeds, condemnation awards, rights to payment of any and every kind and other forms 
 of obligations and receivables, instruments, securities, financial assets and other property which at any time constitute all or part of or are included in the 
 proceeds of any of the T rust Estate (collectively , the “  Collateral  ”). 
 The Indenture T rustee, on behalf of the Noteholders, acknowledges the foregoing Grant, accepts the trusts under this Indenture and 
 agrees to perform its duties required in this Indenture in accordance with the provisions of this Indenture. 
 The foregoing Grant is made in trust to secure (i) the payment of principal of and interest on, and any other amounts owing in respect of, 
 the Notes, equally and ratably without prejudice, priority or distinction except as set forth herein and (ii) compliance with the provisions of this 
 Indenture, all as provided in this Indenture. 
 Without limiting the foregoing Grant, any Receivable repurchased or purchased by (a) the Servicer pursuant to  Section 3.6  of the Servicing 
 Agreement or (b) by the Bank pursuant to  Section 3.4  of the Purchase Agreement shall be deemed to be automatically released from the lien of this 
 Indenture without any action being taken by the Indenture T rustee upon payment by the applicable purchaser of the related Repurchase Price for such 
 Repurchased Receivable. 
 COPAR 2023-1 Indentur e 
--- [START OF PAGE 12] ---
 ARTICLE I 
 DEFINITIONS AND INCORPORA TION BY REFERENCE 
 SECTION 1.1  Definitions  . Capitalized terms are used  in this Indenture as defined in  A  p  p  endix A  to the  Sale Agreement, dated as of the 
 date hereof (as amended, supplemented, or otherwise modified and in ef fect from time to time, the “  Sale  Agreement  ”), between the Issuer and Capital 
 One Auto Receivables, LLC, which also contains rules as to usage that are applicable herein. 
 SECTION 1.2  Incorporation by Reference of T rust Indenture  Act  . Whenever this Indenture refers to a provision  of the TIA, the provision 
 is incorporated by reference in and made a part of this Indenture. The following TIA terms used in this Indenture have the following meanings: 
 “  Commission  ” means the Securities and Exchange Commission. 
 “  indenture securities  ” means the Notes. 
 “  indenture securi  ty  holder  ” means a Noteholder . 
 “  indenture to be qualified  ” means this Indenture. 
 “  indenture trustee  ” or “  institutional trustee  ” means  the Indenture T rustee. 
 “  obl  ig  or  ” on the indenture securities means the Issuer  and any other obligor on the indenture securities. 
 All other TIA terms used in this Indenture that are defined by the TIA, defined by TIA reference to another statute or defined by 
 Commission rule have the meaning assigned to them by such definitions. 
 SECTION 1.3  Other Interpretive Provisions  . All terms  defined in this Indenture shall have the defined meanings when used in any 
 certificate or other document delivered pursuant hereto unless otherwise defined therein. For purposes of this Indenture and all such certificates and 
 other documents, unless the context otherwise requires: (a) accounting terms not otherwise defined in this Indenture, and accounting terms partly 
 defined in this Indenture to the extent not defined, shall have the respective meanings given to them under GAAP (  provided  , that, to the extent that the 
 definitions in this Indenture and GAAP conflict, the definitions in this Indenture shall control); (b) terms defined in Article 9 of the UCC as in ef fect in 
 the relevant jurisdiction and not otherwise defined in this Indenture are used as defined in that Article; (c) the words “hereof,” “herein” and “hereunder” 
 and words of similar import refer to this Indenture as a whole and not to any particular provision of this Indenture; (d) references to any Article, 
 Section, Schedule, Appendix or Exhibit are references to Articles, Sections, Schedules, Appendices and Exhibits in or to this Indenture and references 
 to any paragraph, subsection, clause or other subdivision within any Section or definition refer to such paragraph, subsection, clause or other 
--- CONTRACT TEXT ENDS ---
'@
$ws.Range("E12").Value = $v_E12

$v_F12 = @'
{
  "errors": [
    {
      "location": "Page 3, Section 2.11",
      "error": "The capitalized term 'Clearing Agency' is used but not defined within this document. While Section 1.1 points to an external document for definitions, as a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Clearing Agency' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 2.11 Notices to Clearing Agency"
    },
    {
      "location": "Page 3, Section 2.12",
      "error": "The capitalized term 'Definitive Notes' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Definitive Notes' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 2.12 Definitive Notes"
    },
    {
      "location": "Page 3, Section 2.13",
      "error": "The capitalized term 'Authenticating Agents' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Authenticating Agents' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 2.13 Authenticating Agents"
    },
    {
      "location": "Page 3, Section 2.14",
      "error": "The capitalized term 'Paying Agent' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Paying Agent' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 2.14 Paying Agent"
    },
    {
      "location": "Page 3, Section 2.18",
      "error": "The capitalized term 'Certificate' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Certificate' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 2.18 Transfer Restrictions on Certain Notes Upon a Sale of a Certificate"
    },
    {
      "location": "Page 3, Section 2.19",
      "error": "The capitalized term '144A Notes' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for '144A Notes' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 2.19 Certain Transfer Restrictions on the 144A Notes"
    },
    {
      "location": "Page 5, Section 3.7",
      "error": "The capitalized term 'Receivables' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Receivables' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 3.7 Performance of Obligations; Servicing of Receivables"
    },
    {
      "location": "Page 5, Section 3.11",
      "error": "The capitalized term 'Restricted Payments' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Restricted Payments' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 3.11 Restricted Payments"
    },
    {
      "location": "Page 5, Section 3.12",
      "error": "The capitalized term 'Events of Default' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Events of Default' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 3.12 Notice of Events of Default"
    },
    {
      "location": "Page 5, Section 3.15",
      "error": "The capitalized term 'Administrator' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Administrator' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 3.15 Removal of Administrator"
    },
    {
      "location": "Page 5, Section 4.2",
      "error": "The capitalized term 'Trust Money' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Trust Money' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 4.2 Application of Trust Money"
    },
    {
      "location": "Page 5, Section 5.3",
      "error": "The capitalized term 'Indebtedness' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Indebtedness' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 5.3 Collection of Indebtedness and Suits for Enforcement by the Indenture Trustee"
    },
    {
      "location": "Page 5, Section 5.7",
      "error": "The capitalized term 'Noteholders' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Noteholders' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 5.7 Rights of Noteholders to Receive Principal and Interest"
    },
    {
      "location": "Page 5, Section 5.12",
      "error": "The capitalized term 'Past Defaults' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Past Defaults' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 5.12 Waiver of Past Defaults"
    },
    {
      "location": "Page 7, Section 6.9",
      "error": "The capitalized term 'Merger' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Merger' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 6.9 Successor Indenture Trustee by Mer ger"
    },
    {
      "location": "Page 7, Section 6.10",
      "error": "The capitalized term 'Co-Indenture Trustee' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Co-Indenture Trustee' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 6.10 Appointment of Co-Indenture T rustee or Separate Indenture T rustee"
    },
    {
      "location": "Page 7, Section 6.10",
      "error": "The capitalized term 'Separate Indenture Trustee' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Separate Indenture Trustee' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 6.10 Appointment of Co-Indenture T rustee or Separate Indenture T rustee"
    },
    {
      "location": "Page 7, Section 7.4",
      "error": "The capitalized term 'Certificateholders' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Certificateholders' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 7.4 Statements to Certificateholders and Noteholders"
    },
    {
      "location": "Page 7, Section 7.5",
      "error": "The capitalized term 'Repurchase' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Repurchase' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 7.5 Noteholder Demand for Repurchase, Dispute Resolution"
    },
    {
      "location": "Page 7, Section 7.6",
      "error": "The capitalized term 'Investor' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Investor' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 7.6 Investor Action to Initiate an Asset Review"
    },
    {
      "location": "Page 7, Section 7.6",
      "error": "The capitalized term 'Asset Review' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Asset Review' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 7.6 Investor Action to Initiate an Asset Review"
    },
    {
      "location": "Page 7, Section 8.2",
      "error": "The capitalized term 'Trust Accounts' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Trust Accounts' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 8.2 Trust Accounts"
    },
    {
      "location": "Page 8, Section 8.5",
      "error": "The capitalized term 'Distributions' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Distributions' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 8.5 Distributions"
    },
    {
      "location": "Page 8, Section 10.1",
      "error": "The capitalized term 'Redemption' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Redemption' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 10.1 Redemption"
    },
    {
      "location": "Page 8, Section 10.2",
      "error": "The capitalized term 'Redemption Notice' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Redemption Notice' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 10.2 Form of Redemption Notice"
    },
    {
      "location": "Page 8, Section 10.3",
      "error": "The capitalized term 'Redemption Date' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Redemption Date' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 10.3 Notes Payable on Redemption Date"
    },
    {
      "location": "Page 8, Section 11.13",
      "error": "The capitalized term 'Legal Holidays' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Legal Holidays' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 11.13 Legal Holidays"
    },
    {
      "location": "Page 10, Section 11.14",
      "error": "The capitalized term 'GOVERNING LAW' is used but not defined within this document. Additionally, section titles in the Table of Contents should typically be Title Case for consistency (e.g., 'Governing Law').",
      "suggestion": "Either provide a definition for 'Governing Law' in this document or explicitly state its definition by reference within the relevant section. Also, revise the capitalization of the section title in the Table of Contents to 'Governing Law'.",
      "exact_quote": "Section 11.14 GOVERNING LA W"
    },
    {
      "location": "Page 10, Section 11.19",
      "error": "The capitalized term 'Jury Trial' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Jury Trial' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 11.19 Submission to Jurisdiction; W aiver of Jury T rial"
    },
    {
      "location": "Page 10, Section 11.21",
      "error": "The capitalized term 'U.S.A Patriot Act' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'U.S.A Patriot Act' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 11.21 U.S.A Patriot Act"
    },
    {
      "location": "Page 10, Section 12.2",
      "error": "The capitalized term 'FDIC Rule' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'FDIC Rule' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 12.2 Requirements of the FDIC Rule"
    },
    {
      "location": "Page 10, Section 12.4",
      "error": "The capitalized term 'Repudiation' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Repudiation' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 12.4 Actions Upon Repudiation"
    },
    {
      "location": "Page 11, Introductory Paragraph",
      "error": "The capitalized term 'Holders' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition. It should ideally be consistently 'Noteholders' or explicitly defined.",
      "suggestion": "Either define 'Holders' or consistently use 'Noteholders' throughout the document. If used interchangeably, a definition confirming this should be provided.",
      "exact_quote": "for the benefit of the other party and the equal and ratable benefit of the Holders of the Issuer’s Class A-1"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The word 'Grants' is incorrectly capitalized as a verb. It should be lowercase.",
      "suggestion": "Change 'Grants' to 'grants'.",
      "exact_quote": "hereby Grants in trust"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Closing Date' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Closing Date' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "to the Indenture T rustee on the Closing Date"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Trust Estate' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Trust Estate' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "in and to (i) the T rust Estate"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Servicer' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Servicer' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "by (a) the Servicer pursuant to Section 3.6 of the Servicing Agreement"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Servicing Agreement' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Servicing Agreement' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 3.6 of the Servicing Agreement"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Bank' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Bank' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "or (b) by the Bank pursuant to Section 3.4 of the Purchase Agreement"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Purchase Agreement' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Purchase Agreement' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "Section 3.4 of the Purchase Agreement"
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Repurchase Price' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Repurchase Price' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "of the related Repurchase Price for such Repurchased Receivable."
    },
    {
      "location": "Page 11, Granting Clause",
      "error": "The capitalized term 'Repurchased Receivable' is used but not defined within this document. As a standalone document fragment, a reader cannot access the definition.",
      "suggestion": "Either provide a definition for 'Repurchased Receivable' in this document or explicitly state its definition by reference within the relevant section.",
      "exact_quote": "of the related Repurchase Price for such Repurchased Receivable."
    },
    {
      "location": "Page 12, Section 1.3",
      "error": "The capitalized term 'UCC' is used but not defined within this document. While Section 1.3 points to Article 9 of the UCC for definitions, as a standalone document fragment, a reader cannot access the definition or confirm what 'UCC' refers to without external knowledge.",
      "suggestion": "Either provide a definition for 'UCC' (e.g., Uniform Commercial Code) in this document or explicitly define its meaning by reference within the relevant section.",
      "exact_quote": "Article 9 of the UCC as in ef fect in the relevant jurisdiction"
    }
  ]
}
'@
$ws.Range("F12").Value = $v_F12
